# Electronics catalog refresh (2022-09-10 20:58:20 crawl):
#  - product list re-ordered/re-scraped; one discontinued item
#    (Tefal Dampfbuegeleisen Express Protect SV9202S0, id 6821480) dropped
#  - every remaining row gets the new crawl timestamp
#  - sheet dimension shrinks from A1:O94 to A1:O93
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: id 4090845
$ws.Cells.Item(2,1).NumberFormat = "@"
$ws.Cells.Item(2,1).Value = '4090845'
$ws.Cells.Item(2,2).Value = 'Prix Garantie Batterien AA/LR6 10 Stück'
$ws.Cells.Item(2,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/prix-garantie-batterien-aalr6-10-stueck/p/4090845'
$ws.Cells.Item(2,4).Value = '10ST'
$ws.Cells.Item(2,7).Value = 'Coop'
$ws.Cells.Item(2,8).NumberFormat = "@"
$ws.Cells.Item(2,8).Value = '4.75'
$ws.Cells.Item(2,9).Value = '0.48/1ST'
$ws.Cells.Item(2,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(2,11).NumberFormat = "@"
$ws.Cells.Item(2,11).Value = '0.48'
$ws.Cells.Item(2,12).Value = '1ST'
$ws.Cells.Item(2,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Cells.Item(2,14).Value = 'Prix Garantie Batterien AA/LR6 10 Stück 4.75 Schweizer Franken'
$ws.Cells.Item(2,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(2,5).Value = 11
$ws.Cells.Item(2,6).Value = 3
# Row 3: id 4090846
$ws.Cells.Item(3,1).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = '4090846'
$ws.Cells.Item(3,2).Value = 'Prix Garantie Batterien AAA/LR03 10 Stück'
$ws.Cells.Item(3,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/prix-garantie-batterien-aaalr03-10-stueck/p/4090846'
$ws.Cells.Item(3,4).Value = '10ST'
$ws.Cells.Item(3,7).Value = 'Coop'
$ws.Cells.Item(3,8).NumberFormat = "@"
$ws.Cells.Item(3,8).Value = '4.70'
$ws.Cells.Item(3,9).Value = '0.47/1ST'
$ws.Cells.Item(3,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(3,11).NumberFormat = "@"
$ws.Cells.Item(3,11).Value = '0.47'
$ws.Cells.Item(3,12).Value = '1ST'
$ws.Cells.Item(3,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Cells.Item(3,14).Value = 'Prix Garantie Batterien AAA/LR03 10 Stück 4.70 Schweizer Franken'
$ws.Cells.Item(3,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(3,5).Value = 13
$ws.Cells.Item(3,6).Value = 4
# Row 4: id 5920020
$ws.Cells.Item(4,1).NumberFormat = "@"
$ws.Cells.Item(4,1).Value = '5920020'
$ws.Cells.Item(4,2).Value = 'BASIC Alkaline Batterien LR6/AA, 12 Stück'
$ws.Cells.Item(4,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/basic-alkaline-batterien-lr6aa-12-stueck/p/5920020'
$ws.Cells.Item(4,4).Value = '12ST'
$ws.Cells.Item(4,7).Value = 'Coop'
$ws.Cells.Item(4,8).NumberFormat = "@"
$ws.Cells.Item(4,8).Value = '9.95'
$ws.Cells.Item(4,9).Value = '0.83/1ST'
$ws.Cells.Item(4,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(4,11).NumberFormat = "@"
$ws.Cells.Item(4,11).Value = '0.83'
$ws.Cells.Item(4,12).Value = '1ST'
$ws.Cells.Item(4,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Cells.Item(4,14).Value = 'BASIC Alkaline Batterien LR6/AA, 12 Stück 9.95 Schweizer Franken'
$ws.Cells.Item(4,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(4,5).ClearContents()
$ws.Cells.Item(4,6).Value = 0
# Row 5: id 5920351
$ws.Cells.Item(5,1).NumberFormat = "@"
$ws.Cells.Item(5,1).Value = '5920351'
$ws.Cells.Item(5,2).Value = 'BASIC Alkaline Batterien LR03/AAA 12 Stück'
$ws.Cells.Item(5,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/basic-alkaline-batterien-lr03aaa-12-stueck/p/5920351'
$ws.Cells.Item(5,4).Value = '12ST'
$ws.Cells.Item(5,7).Value = 'Coop'
$ws.Cells.Item(5,8).NumberFormat = "@"
$ws.Cells.Item(5,8).Value = '9.95'
$ws.Cells.Item(5,9).Value = '0.83/1ST'
$ws.Cells.Item(5,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(5,11).NumberFormat = "@"
$ws.Cells.Item(5,11).Value = '0.83'
$ws.Cells.Item(5,12).Value = '1ST'
$ws.Cells.Item(5,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Cells.Item(5,14).Value = 'BASIC Alkaline Batterien LR03/AAA 12 Stück 9.95 Schweizer Franken'
$ws.Cells.Item(5,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 4
# Row 6: id 4488074
$ws.Cells.Item(6,1).NumberFormat = "@"
$ws.Cells.Item(6,1).Value = '4488074'
$ws.Cells.Item(6,2).Value = 'Varta Knopfzellen CR2032 2 Stück'
$ws.Cells.Item(6,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-knopfzellen-cr2032-2-stueck/p/4488074'
$ws.Cells.Item(6,4).Value = '2ST'
$ws.Cells.Item(6,7).Value = 'Varta'
$ws.Cells.Item(6,8).NumberFormat = "@"
$ws.Cells.Item(6,8).Value = '8.95'
$ws.Cells.Item(6,9).Value = '4.48/1ST'
$ws.Cells.Item(6,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(6,11).NumberFormat = "@"
$ws.Cells.Item(6,11).Value = '4.48'
$ws.Cells.Item(6,12).Value = '1ST'
$ws.Cells.Item(6,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(6,14).Value = 'Varta Knopfzellen CR2032 2 Stück 8.95 Schweizer Franken'
$ws.Cells.Item(6,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(6,5).Value = 4
$ws.Cells.Item(6,6).Value = 4.5
# Row 7: id 6762335
$ws.Cells.Item(7,1).NumberFormat = "@"
$ws.Cells.Item(7,1).Value = '6762335'
$ws.Cells.Item(7,2).Value = 'Duracell Batterien Optimum AA/LR6 8 Stück'
$ws.Cells.Item(7,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/duracell-batterien-optimum-aalr6-8-stueck/p/6762335'
$ws.Cells.Item(7,4).Value = '8ST'
$ws.Cells.Item(7,7).Value = 'Duracell'
$ws.Cells.Item(7,8).NumberFormat = "@"
$ws.Cells.Item(7,8).Value = '19.95'
$ws.Cells.Item(7,9).Value = '2.49/1ST'
$ws.Cells.Item(7,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(7,11).NumberFormat = "@"
$ws.Cells.Item(7,11).Value = '2.49'
$ws.Cells.Item(7,12).Value = '1ST'
$ws.Cells.Item(7,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Cells.Item(7,14).Value = 'Duracell Batterien Optimum AA/LR6 8 Stück 19.95 Schweizer Franken'
$ws.Cells.Item(7,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 5
# Row 8: id 6762336
$ws.Cells.Item(8,1).NumberFormat = "@"
$ws.Cells.Item(8,1).Value = '6762336'
$ws.Cells.Item(8,2).Value = 'Duracell Batterien Optimum AAA/LR03 8 Stück'
$ws.Cells.Item(8,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/duracell-batterien-optimum-aaalr03-8-stueck/p/6762336'
$ws.Cells.Item(8,4).Value = '8ST'
$ws.Cells.Item(8,7).Value = 'Duracell'
$ws.Cells.Item(8,8).NumberFormat = "@"
$ws.Cells.Item(8,8).Value = '19.95'
$ws.Cells.Item(8,9).Value = '2.49/1ST'
$ws.Cells.Item(8,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(8,11).NumberFormat = "@"
$ws.Cells.Item(8,11).Value = '2.49'
$ws.Cells.Item(8,12).Value = '1ST'
$ws.Cells.Item(8,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Cells.Item(8,14).Value = 'Duracell Batterien Optimum AAA/LR03 8 Stück 19.95 Schweizer Franken'
$ws.Cells.Item(8,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(8,5).ClearContents()
$ws.Cells.Item(8,6).Value = 0
# Row 9: id 4014528
$ws.Cells.Item(9,1).NumberFormat = "@"
$ws.Cells.Item(9,1).Value = '4014528'
$ws.Cells.Item(9,2).Value = 'Varta Longlife Batterien AAA/LR03 10 Stück'
$ws.Cells.Item(9,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-batterien-aaalr03-10-stueck/p/4014528'
$ws.Cells.Item(9,4).Value = '10ST'
$ws.Cells.Item(9,7).Value = 'Varta'
$ws.Cells.Item(9,8).NumberFormat = "@"
$ws.Cells.Item(9,8).Value = '16.95'
$ws.Cells.Item(9,9).Value = '1.70/1ST'
$ws.Cells.Item(9,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(9,11).NumberFormat = "@"
$ws.Cells.Item(9,11).Value = '1.70'
$ws.Cells.Item(9,12).Value = '1ST'
$ws.Cells.Item(9,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Cells.Item(9,14).Value = 'Varta Longlife Batterien AAA/LR03 10 Stück 16.95 Schweizer Franken'
$ws.Cells.Item(9,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 5
# Row 10: id 4905479
$ws.Cells.Item(10,1).NumberFormat = "@"
$ws.Cells.Item(10,1).Value = '4905479'
$ws.Cells.Item(10,2).Value = 'Coop Alkaline Batterie, LR6/AA, 4 Stk.'
$ws.Cells.Item(10,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/coop-alkaline-batterie-lr6aa-4-stk/p/4905479'
$ws.Cells.Item(10,4).Value = '4ST'
$ws.Cells.Item(10,7).Value = 'Coop'
$ws.Cells.Item(10,8).NumberFormat = "@"
$ws.Cells.Item(10,8).Value = '4.95'
$ws.Cells.Item(10,9).Value = '1.24/1ST'
$ws.Cells.Item(10,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(10,11).NumberFormat = "@"
$ws.Cells.Item(10,11).Value = '1.24'
$ws.Cells.Item(10,12).Value = '1ST'
$ws.Cells.Item(10,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Cells.Item(10,14).Value = 'Coop Alkaline Batterie, LR6/AA, 4 Stk. 4.95 Schweizer Franken'
$ws.Cells.Item(10,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(10,5).Value = 7
$ws.Cells.Item(10,6).Value = 4
# Row 11: id 4905467
$ws.Cells.Item(11,1).NumberFormat = "@"
$ws.Cells.Item(11,1).Value = '4905467'
$ws.Cells.Item(11,2).Value = 'Alkaline Batterie AAA/LR03 4 Stück'
$ws.Cells.Item(11,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/alkaline-batterie-aaalr03-4-stueck/p/4905467'
$ws.Cells.Item(11,4).Value = '4ST'
$ws.Cells.Item(11,7).Value = 'Coop'
$ws.Cells.Item(11,8).NumberFormat = "@"
$ws.Cells.Item(11,8).Value = '4.95'
$ws.Cells.Item(11,9).Value = '1.24/1ST'
$ws.Cells.Item(11,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(11,11).NumberFormat = "@"
$ws.Cells.Item(11,11).Value = '1.24'
$ws.Cells.Item(11,12).Value = '1ST'
$ws.Cells.Item(11,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Cells.Item(11,14).Value = 'Alkaline Batterie AAA/LR03 4 Stück 4.95 Schweizer Franken'
$ws.Cells.Item(11,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(11,5).Value = 4
$ws.Cells.Item(11,6).Value = 4.5
# Row 12: id 7032801
$ws.Cells.Item(12,1).NumberFormat = "@"
$ws.Cells.Item(12,1).Value = '7032801'
$ws.Cells.Item(12,2).Value = 'Duracell Batterien Plus AA/LR6 16 Stück'
$ws.Cells.Item(12,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/duracell-batterien-plus-aalr6-16-stueck/p/7032801'
$ws.Cells.Item(12,4).Value = '16ST'
$ws.Cells.Item(12,7).Value = 'Duracell'
$ws.Cells.Item(12,8).NumberFormat = "@"
$ws.Cells.Item(12,8).Value = '19.90'
$ws.Cells.Item(12,9).Value = '1.24/1ST'
$ws.Cells.Item(12,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(12,11).NumberFormat = "@"
$ws.Cells.Item(12,11).Value = '1.24'
$ws.Cells.Item(12,12).Value = '1ST'
$ws.Cells.Item(12,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Cells.Item(12,14).Value = 'Duracell Batterien Plus AA/LR6 16 Stück 50% Aktion 19.90 Schweizer Franken statt 39.80 Schweizer Franken'
$ws.Cells.Item(12,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(12,5).ClearContents()
$ws.Cells.Item(12,6).Value = 0
# Row 13: id 4905482
$ws.Cells.Item(13,1).NumberFormat = "@"
$ws.Cells.Item(13,1).Value = '4905482'
$ws.Cells.Item(13,2).Value = 'Alkaline Batterie LR14/C 2 Stück'
$ws.Cells.Item(13,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/alkaline-batterie-lr14c-2-stueck/p/4905482'
$ws.Cells.Item(13,4).Value = '2ST'
$ws.Cells.Item(13,7).Value = 'Coop'
$ws.Cells.Item(13,8).NumberFormat = "@"
$ws.Cells.Item(13,8).Value = '4.95'
$ws.Cells.Item(13,9).Value = '2.48/1ST'
$ws.Cells.Item(13,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(13,11).NumberFormat = "@"
$ws.Cells.Item(13,11).Value = '2.48'
$ws.Cells.Item(13,12).Value = '1ST'
$ws.Cells.Item(13,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(13,14).Value = 'Alkaline Batterie LR14/C 2 Stück 4.95 Schweizer Franken'
$ws.Cells.Item(13,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 4.5
# Row 14: id 6153846
$ws.Cells.Item(14,1).NumberFormat = "@"
$ws.Cells.Item(14,1).Value = '6153846'
$ws.Cells.Item(14,2).Value = 'Duracell Knopfzelle CR2032 3V 2 Stück'
$ws.Cells.Item(14,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-knopfzelle-cr2032-3v-2-stueck/p/6153846'
$ws.Cells.Item(14,4).Value = '2ST'
$ws.Cells.Item(14,7).Value = 'Duracell'
$ws.Cells.Item(14,8).NumberFormat = "@"
$ws.Cells.Item(14,8).Value = '9.95'
$ws.Cells.Item(14,9).Value = '4.98/1ST'
$ws.Cells.Item(14,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(14,11).NumberFormat = "@"
$ws.Cells.Item(14,11).Value = '4.98'
$ws.Cells.Item(14,12).Value = '1ST'
$ws.Cells.Item(14,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(14,14).Value = 'Duracell Knopfzelle CR2032 3V 2 Stück 9.95 Schweizer Franken'
$ws.Cells.Item(14,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(14,5).Value = 1
$ws.Cells.Item(14,6).Value = 5
# Row 15: id 7032804
$ws.Cells.Item(15,1).NumberFormat = "@"
$ws.Cells.Item(15,1).Value = '7032804'
$ws.Cells.Item(15,2).Value = 'Duracell Batterien Plus AAA/LR03 16 Stück'
$ws.Cells.Item(15,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/duracell-batterien-plus-aaalr03-16-stueck/p/7032804'
$ws.Cells.Item(15,4).Value = '16ST'
$ws.Cells.Item(15,7).Value = 'Duracell'
$ws.Cells.Item(15,8).NumberFormat = "@"
$ws.Cells.Item(15,8).Value = '19.90'
$ws.Cells.Item(15,9).Value = '1.24/1ST'
$ws.Cells.Item(15,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(15,11).NumberFormat = "@"
$ws.Cells.Item(15,11).Value = '1.24'
$ws.Cells.Item(15,12).Value = '1ST'
$ws.Cells.Item(15,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Cells.Item(15,14).Value = 'Duracell Batterien Plus AAA/LR03 16 Stück 50% Aktion 19.90 Schweizer Franken statt 39.80 Schweizer Franken'
$ws.Cells.Item(15,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(15,5).ClearContents()
$ws.Cells.Item(15,6).Value = 0
# Row 16: id 4905485
$ws.Cells.Item(16,1).NumberFormat = "@"
$ws.Cells.Item(16,1).Value = '4905485'
$ws.Cells.Item(16,2).Value = 'Alkaline Batterie 6LP3146/9V'
$ws.Cells.Item(16,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/alkaline-batterie-6lp31469v/p/4905485'
$ws.Cells.Item(16,4).Value = '1ST'
$ws.Cells.Item(16,7).Value = 'Coop'
$ws.Cells.Item(16,8).NumberFormat = "@"
$ws.Cells.Item(16,8).Value = '4.95'
$ws.Cells.Item(16,9).Value = '4.95/1ST'
$ws.Cells.Item(16,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(16,11).NumberFormat = "@"
$ws.Cells.Item(16,11).Value = '4.95'
$ws.Cells.Item(16,12).Value = '1ST'
$ws.Cells.Item(16,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(16,14).Value = 'Alkaline Batterie 6LP3146/9V 4.95 Schweizer Franken'
$ws.Cells.Item(16,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(16,5).Value = 2
$ws.Cells.Item(16,6).Value = 4.5
# Row 17: id 6761134
$ws.Cells.Item(17,1).NumberFormat = "@"
$ws.Cells.Item(17,1).Value = '6761134'
$ws.Cells.Item(17,2).Value = 'Duracell Batterien PLUS D/LR20 2 Stück'
$ws.Cells.Item(17,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterien-plus-dlr20-2-stueck/p/6761134'
$ws.Cells.Item(17,4).Value = '2ST'
$ws.Cells.Item(17,7).Value = 'Duracell'
$ws.Cells.Item(17,8).NumberFormat = "@"
$ws.Cells.Item(17,8).Value = '9.95'
$ws.Cells.Item(17,9).Value = '4.98/1ST'
$ws.Cells.Item(17,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(17,11).NumberFormat = "@"
$ws.Cells.Item(17,11).Value = '4.98'
$ws.Cells.Item(17,12).Value = '1ST'
$ws.Cells.Item(17,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(17,14).Value = 'Duracell Batterien PLUS D/LR20 2 Stück 9.95 Schweizer Franken'
$ws.Cells.Item(17,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(17,5).ClearContents()
$ws.Cells.Item(17,6).Value = 0
# Row 18: id 4949707
$ws.Cells.Item(18,1).NumberFormat = "@"
$ws.Cells.Item(18,1).Value = '4949707'
$ws.Cells.Item(18,2).Value = 'Varta Knopfzellen CR2025 2 Stück'
$ws.Cells.Item(18,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-knopfzellen-cr2025-2-stueck/p/4949707'
$ws.Cells.Item(18,4).Value = '2ST'
$ws.Cells.Item(18,7).Value = 'Varta'
$ws.Cells.Item(18,8).NumberFormat = "@"
$ws.Cells.Item(18,8).Value = '8.95'
$ws.Cells.Item(18,9).Value = '4.48/1ST'
$ws.Cells.Item(18,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(18,11).NumberFormat = "@"
$ws.Cells.Item(18,11).Value = '4.48'
$ws.Cells.Item(18,12).Value = '1ST'
$ws.Cells.Item(18,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(18,14).Value = 'Varta Knopfzellen CR2025 2 Stück 8.95 Schweizer Franken'
$ws.Cells.Item(18,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(18,5).ClearContents()
$ws.Cells.Item(18,6).Value = 0
# Row 19: id 6761133
$ws.Cells.Item(19,1).NumberFormat = "@"
$ws.Cells.Item(19,1).Value = '6761133'
$ws.Cells.Item(19,2).Value = 'Duracell Batterien PLUS C/LR14 2 Stück'
$ws.Cells.Item(19,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterien-plus-clr14-2-stueck/p/6761133'
$ws.Cells.Item(19,4).Value = '2ST'
$ws.Cells.Item(19,7).Value = 'Duracell'
$ws.Cells.Item(19,8).NumberFormat = "@"
$ws.Cells.Item(19,8).Value = '9.95'
$ws.Cells.Item(19,9).Value = '4.98/1ST'
$ws.Cells.Item(19,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(19,11).NumberFormat = "@"
$ws.Cells.Item(19,11).Value = '4.98'
$ws.Cells.Item(19,12).Value = '1ST'
$ws.Cells.Item(19,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(19,14).Value = 'Duracell Batterien PLUS C/LR14 2 Stück 9.95 Schweizer Franken'
$ws.Cells.Item(19,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(19,5).ClearContents()
$ws.Cells.Item(19,6).Value = 0
# Row 20: id 4014527
$ws.Cells.Item(20,1).NumberFormat = "@"
$ws.Cells.Item(20,1).Value = '4014527'
$ws.Cells.Item(20,2).Value = 'Varta Longlife Batterien AA/LR6 10 Stück'
$ws.Cells.Item(20,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-batterien-aalr6-10-stueck/p/4014527'
$ws.Cells.Item(20,4).Value = '10ST'
$ws.Cells.Item(20,7).Value = 'Varta'
$ws.Cells.Item(20,8).NumberFormat = "@"
$ws.Cells.Item(20,8).Value = '15.95'
$ws.Cells.Item(20,9).Value = '1.60/1ST'
$ws.Cells.Item(20,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(20,11).NumberFormat = "@"
$ws.Cells.Item(20,11).Value = '1.60'
$ws.Cells.Item(20,12).Value = '1ST'
$ws.Cells.Item(20,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Cells.Item(20,14).Value = 'Varta Longlife Batterien AA/LR6 10 Stück 15.95 Schweizer Franken'
$ws.Cells.Item(20,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(20,5).Value = 4
$ws.Cells.Item(20,6).Value = 3
# Row 21: id 6753556
$ws.Cells.Item(21,1).NumberFormat = "@"
$ws.Cells.Item(21,1).Value = '6753556'
$ws.Cells.Item(21,2).Value = 'Duracell Batterien PLUS AA/LR6 4 Stück'
$ws.Cells.Item(21,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/duracell-batterien-plus-aalr6-4-stueck/p/6753556'
$ws.Cells.Item(21,4).Value = '4ST'
$ws.Cells.Item(21,7).Value = 'Duracell'
$ws.Cells.Item(21,8).NumberFormat = "@"
$ws.Cells.Item(21,8).Value = '9.95'
$ws.Cells.Item(21,9).Value = '2.49/1ST'
$ws.Cells.Item(21,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(21,11).NumberFormat = "@"
$ws.Cells.Item(21,11).Value = '2.49'
$ws.Cells.Item(21,12).Value = '1ST'
$ws.Cells.Item(21,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Cells.Item(21,14).Value = 'Duracell Batterien PLUS AA/LR6 4 Stück 9.95 Schweizer Franken'
$ws.Cells.Item(21,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(21,5).ClearContents()
$ws.Cells.Item(21,6).Value = 0
# Row 22: id 3494066
$ws.Cells.Item(22,1).NumberFormat = "@"
$ws.Cells.Item(22,1).Value = '3494066'
$ws.Cells.Item(22,2).Value = 'Varta Longlife Power 9V 1er Bli'
$ws.Cells.Item(22,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-power-9v-1er-bli/p/3494066'
$ws.Cells.Item(22,4).Value = '1ST'
$ws.Cells.Item(22,7).Value = 'Varta'
$ws.Cells.Item(22,8).NumberFormat = "@"
$ws.Cells.Item(22,8).Value = '8.95'
$ws.Cells.Item(22,9).Value = '8.95/1ST'
$ws.Cells.Item(22,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(22,11).NumberFormat = "@"
$ws.Cells.Item(22,11).Value = '8.95'
$ws.Cells.Item(22,12).Value = '1ST'
$ws.Cells.Item(22,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(22,14).Value = 'Varta Longlife Power 9V 1er Bli 8.95 Schweizer Franken'
$ws.Cells.Item(22,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(22,5).Value = 2
$ws.Cells.Item(22,6).Value = 4.5
# Row 23: id 3494865
$ws.Cells.Item(23,1).NumberFormat = "@"
$ws.Cells.Item(23,1).Value = '3494865'
$ws.Cells.Item(23,2).Value = 'Varta Electronics CR2430 1er Bli'
$ws.Cells.Item(23,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-cr2430-1er-bli/p/3494865'
$ws.Cells.Item(23,4).Value = '1ST'
$ws.Cells.Item(23,7).Value = 'Varta'
$ws.Cells.Item(23,8).NumberFormat = "@"
$ws.Cells.Item(23,8).Value = '4.95'
$ws.Cells.Item(23,9).Value = '4.95/1ST'
$ws.Cells.Item(23,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(23,11).NumberFormat = "@"
$ws.Cells.Item(23,11).Value = '4.95'
$ws.Cells.Item(23,12).Value = '1ST'
$ws.Cells.Item(23,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(23,14).Value = 'Varta Electronics CR2430 1er Bli 4.95 Schweizer Franken'
$ws.Cells.Item(23,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(23,5).Value = 1
$ws.Cells.Item(23,6).Value = 5
# Row 24: id 6577801
$ws.Cells.Item(24,1).NumberFormat = "@"
$ws.Cells.Item(24,1).Value = '6577801'
$ws.Cells.Item(24,2).Value = 'Duracell Knopfzelle CR2025 2 Stück'
$ws.Cells.Item(24,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-knopfzelle-cr2025-2-stueck/p/6577801'
$ws.Cells.Item(24,4).Value = '2ST'
$ws.Cells.Item(24,7).Value = 'Duracell'
$ws.Cells.Item(24,8).NumberFormat = "@"
$ws.Cells.Item(24,8).Value = '9.95'
$ws.Cells.Item(24,9).Value = '4.98/1ST'
$ws.Cells.Item(24,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(24,11).NumberFormat = "@"
$ws.Cells.Item(24,11).Value = '4.98'
$ws.Cells.Item(24,12).Value = '1ST'
$ws.Cells.Item(24,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(24,14).Value = 'Duracell Knopfzelle CR2025 2 Stück 9.95 Schweizer Franken'
$ws.Cells.Item(24,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(24,5).ClearContents()
$ws.Cells.Item(24,6).Value = 0
# Row 25: id 6761135
$ws.Cells.Item(25,1).NumberFormat = "@"
$ws.Cells.Item(25,1).Value = '6761135'
$ws.Cells.Item(25,2).Value = 'Duracell Batterie PLUS 9V/6LR61 1 Stück'
$ws.Cells.Item(25,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterie-plus-9v6lr61-1-stueck/p/6761135'
$ws.Cells.Item(25,4).Value = '1ST'
$ws.Cells.Item(25,7).Value = 'Duracell'
$ws.Cells.Item(25,8).NumberFormat = "@"
$ws.Cells.Item(25,8).Value = '9.95'
$ws.Cells.Item(25,9).Value = '9.95/1ST'
$ws.Cells.Item(25,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(25,11).NumberFormat = "@"
$ws.Cells.Item(25,11).Value = '9.95'
$ws.Cells.Item(25,12).Value = '1ST'
$ws.Cells.Item(25,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(25,14).Value = 'Duracell Batterie PLUS 9V/6LR61 1 Stück 9.95 Schweizer Franken'
$ws.Cells.Item(25,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(25,5).ClearContents()
$ws.Cells.Item(25,6).Value = 0
# Row 26: id 3494233
$ws.Cells.Item(26,1).NumberFormat = "@"
$ws.Cells.Item(26,1).Value = '3494233'
$ws.Cells.Item(26,2).Value = 'Varta Electronics CR2032 1er Bli'
$ws.Cells.Item(26,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-cr2032-1er-bli/p/3494233'
$ws.Cells.Item(26,4).Value = '1ST'
$ws.Cells.Item(26,7).Value = 'Varta'
$ws.Cells.Item(26,8).NumberFormat = "@"
$ws.Cells.Item(26,8).Value = '4.95'
$ws.Cells.Item(26,9).Value = '4.95/1ST'
$ws.Cells.Item(26,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(26,11).NumberFormat = "@"
$ws.Cells.Item(26,11).Value = '4.95'
$ws.Cells.Item(26,12).Value = '1ST'
$ws.Cells.Item(26,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(26,14).Value = 'Varta Electronics CR2032 1er Bli 4.95 Schweizer Franken'
$ws.Cells.Item(26,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(26,5).Value = 2
$ws.Cells.Item(26,6).Value = 4.5
# Row 27: id 3494234
$ws.Cells.Item(27,1).NumberFormat = "@"
$ws.Cells.Item(27,1).Value = '3494234'
$ws.Cells.Item(27,2).Value = 'Varta Electronics V13GA 1er Bli'
$ws.Cells.Item(27,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-v13ga-1er-bli/p/3494234'
$ws.Cells.Item(27,4).Value = '1ST'
$ws.Cells.Item(27,7).Value = 'Varta'
$ws.Cells.Item(27,8).NumberFormat = "@"
$ws.Cells.Item(27,8).Value = '2.95'
$ws.Cells.Item(27,9).Value = '2.95/1ST'
$ws.Cells.Item(27,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(27,11).NumberFormat = "@"
$ws.Cells.Item(27,11).Value = '2.95'
$ws.Cells.Item(27,12).Value = '1ST'
$ws.Cells.Item(27,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(27,14).Value = 'Varta Electronics V13GA 1er Bli 2.95 Schweizer Franken'
$ws.Cells.Item(27,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(27,5).ClearContents()
$ws.Cells.Item(27,6).Value = 0
# Row 28: id 6753557
$ws.Cells.Item(28,1).NumberFormat = "@"
$ws.Cells.Item(28,1).Value = '6753557'
$ws.Cells.Item(28,2).Value = 'Duracell Batterien PLUS AAA/LR03 4 Stück'
$ws.Cells.Item(28,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/duracell-batterien-plus-aaalr03-4-stueck/p/6753557'
$ws.Cells.Item(28,4).Value = '4ST'
$ws.Cells.Item(28,7).Value = 'Duracell'
$ws.Cells.Item(28,8).NumberFormat = "@"
$ws.Cells.Item(28,8).Value = '9.95'
$ws.Cells.Item(28,9).Value = '2.49/1ST'
$ws.Cells.Item(28,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(28,11).NumberFormat = "@"
$ws.Cells.Item(28,11).Value = '2.49'
$ws.Cells.Item(28,12).Value = '1ST'
$ws.Cells.Item(28,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Cells.Item(28,14).Value = 'Duracell Batterien PLUS AAA/LR03 4 Stück 9.95 Schweizer Franken'
$ws.Cells.Item(28,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(28,5).ClearContents()
$ws.Cells.Item(28,6).Value = 0
# Row 29: id 4145012
$ws.Cells.Item(29,1).NumberFormat = "@"
$ws.Cells.Item(29,1).Value = '4145012'
$ws.Cells.Item(29,2).Value = 'LED Teelichter warmweiss 4 Stück'
$ws.Cells.Item(29,3).Value = '/de/haushalt-tier/elektroartikel-batterien/beleuchtung/taschenlampen-lichter/led-teelichter-warmweiss-4-stueck/p/4145012'
$ws.Cells.Item(29,4).ClearContents()
$ws.Cells.Item(29,7).Value = 'Coop'
$ws.Cells.Item(29,8).NumberFormat = "@"
$ws.Cells.Item(29,8).Value = '8.95'
$ws.Cells.Item(29,9).ClearContents()
$ws.Cells.Item(29,10).ClearContents()
$ws.Cells.Item(29,11).ClearContents()
$ws.Cells.Item(29,12).ClearContents()
$ws.Cells.Item(29,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''beleuchtung'', ''taschenlampen-lichter'']'
$ws.Cells.Item(29,14).Value = 'LED Teelichter warmweiss 4 Stück 8.95 Schweizer Franken'
$ws.Cells.Item(29,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(29,5).ClearContents()
$ws.Cells.Item(29,6).Value = 0
# Row 30: id 6753554
$ws.Cells.Item(30,1).NumberFormat = "@"
$ws.Cells.Item(30,1).Value = '6753554'
$ws.Cells.Item(30,2).Value = 'Duracell Batterien Optimum AA/LR6 4 Stück'
$ws.Cells.Item(30,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/duracell-batterien-optimum-aalr6-4-stueck/p/6753554'
$ws.Cells.Item(30,4).Value = '4ST'
$ws.Cells.Item(30,7).Value = 'Duracell'
$ws.Cells.Item(30,8).NumberFormat = "@"
$ws.Cells.Item(30,8).Value = '11.95'
$ws.Cells.Item(30,9).Value = '2.99/1ST'
$ws.Cells.Item(30,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(30,11).NumberFormat = "@"
$ws.Cells.Item(30,11).Value = '2.99'
$ws.Cells.Item(30,12).Value = '1ST'
$ws.Cells.Item(30,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Cells.Item(30,14).Value = 'Duracell Batterien Optimum AA/LR6 4 Stück 11.95 Schweizer Franken'
$ws.Cells.Item(30,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(30,5).Value = 1
$ws.Cells.Item(30,6).Value = 5
# Row 31: id 4905484
$ws.Cells.Item(31,1).NumberFormat = "@"
$ws.Cells.Item(31,1).Value = '4905484'
$ws.Cells.Item(31,2).Value = 'Alkaline Batterie LR20/D 2 Stück'
$ws.Cells.Item(31,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/alkaline-batterie-lr20d-2-stueck/p/4905484'
$ws.Cells.Item(31,4).Value = '2ST'
$ws.Cells.Item(31,7).Value = 'Coop'
$ws.Cells.Item(31,8).NumberFormat = "@"
$ws.Cells.Item(31,8).Value = '5.95'
$ws.Cells.Item(31,9).Value = '2.98/1ST'
$ws.Cells.Item(31,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(31,11).NumberFormat = "@"
$ws.Cells.Item(31,11).Value = '2.98'
$ws.Cells.Item(31,12).Value = '1ST'
$ws.Cells.Item(31,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(31,14).Value = 'Alkaline Batterie LR20/D 2 Stück 5.95 Schweizer Franken'
$ws.Cells.Item(31,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(31,5).Value = 3
$ws.Cells.Item(31,6).Value = 5
# Row 32: id 3494133
$ws.Cells.Item(32,1).NumberFormat = "@"
$ws.Cells.Item(32,1).Value = '3494133'
$ws.Cells.Item(32,2).Value = 'Varta Longlife Power AA 8er Bli'
$ws.Cells.Item(32,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-power-aa-8er-bli/p/3494133'
$ws.Cells.Item(32,4).Value = '8ST'
$ws.Cells.Item(32,7).Value = 'Varta'
$ws.Cells.Item(32,8).NumberFormat = "@"
$ws.Cells.Item(32,8).Value = '14.95'
$ws.Cells.Item(32,9).Value = '1.87/1ST'
$ws.Cells.Item(32,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(32,11).NumberFormat = "@"
$ws.Cells.Item(32,11).Value = '1.87'
$ws.Cells.Item(32,12).Value = '1ST'
$ws.Cells.Item(32,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Cells.Item(32,14).Value = 'Varta Longlife Power AA 8er Bli 14.95 Schweizer Franken'
$ws.Cells.Item(32,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(32,5).Value = 3
$ws.Cells.Item(32,6).Value = 3
# Row 33: id 5683906
$ws.Cells.Item(33,1).NumberFormat = "@"
$ws.Cells.Item(33,1).Value = '5683906'
$ws.Cells.Item(33,2).Value = 'Prix Garantie Wasserkocher 1l'
$ws.Cells.Item(33,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/prix-garantie-wasserkocher-1l/p/5683906'
$ws.Cells.Item(33,4).Value = '1ST'
$ws.Cells.Item(33,7).Value = 'Coop'
$ws.Cells.Item(33,8).NumberFormat = "@"
$ws.Cells.Item(33,8).Value = '19.95'
$ws.Cells.Item(33,9).Value = '19.95/1ST'
$ws.Cells.Item(33,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(33,11).NumberFormat = "@"
$ws.Cells.Item(33,11).Value = '19.95'
$ws.Cells.Item(33,12).Value = '1ST'
$ws.Cells.Item(33,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''kuechengeraete'']'
$ws.Cells.Item(33,14).Value = 'Prix Garantie Wasserkocher 1l 19.95 Schweizer Franken'
$ws.Cells.Item(33,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(33,5).Value = 5
$ws.Cells.Item(33,6).Value = 4
# Row 34: id 3494130
$ws.Cells.Item(34,1).NumberFormat = "@"
$ws.Cells.Item(34,1).Value = '3494130'
$ws.Cells.Item(34,2).Value = 'Varta Longlife Power AA 4er Bli'
$ws.Cells.Item(34,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-power-aa-4er-bli/p/3494130'
$ws.Cells.Item(34,4).Value = '4ST'
$ws.Cells.Item(34,7).Value = 'Varta'
$ws.Cells.Item(34,8).NumberFormat = "@"
$ws.Cells.Item(34,8).Value = '8.95'
$ws.Cells.Item(34,9).Value = '2.24/1ST'
$ws.Cells.Item(34,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(34,11).NumberFormat = "@"
$ws.Cells.Item(34,11).Value = '2.24'
$ws.Cells.Item(34,12).Value = '1ST'
$ws.Cells.Item(34,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Cells.Item(34,14).Value = 'Varta Longlife Power AA 4er Bli 8.95 Schweizer Franken'
$ws.Cells.Item(34,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(34,5).Value = 1
$ws.Cells.Item(34,6).Value = 4
# Row 35: id 3494131
$ws.Cells.Item(35,1).NumberFormat = "@"
$ws.Cells.Item(35,1).Value = '3494131'
$ws.Cells.Item(35,2).Value = 'Varta Longlife Power AAA 4er Bli'
$ws.Cells.Item(35,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-power-aaa-4er-bli/p/3494131'
$ws.Cells.Item(35,4).Value = '4ST'
$ws.Cells.Item(35,7).Value = 'Varta'
$ws.Cells.Item(35,8).NumberFormat = "@"
$ws.Cells.Item(35,8).Value = '8.95'
$ws.Cells.Item(35,9).Value = '2.24/1ST'
$ws.Cells.Item(35,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(35,11).NumberFormat = "@"
$ws.Cells.Item(35,11).Value = '2.24'
$ws.Cells.Item(35,12).Value = '1ST'
$ws.Cells.Item(35,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Cells.Item(35,14).Value = 'Varta Longlife Power AAA 4er Bli 8.95 Schweizer Franken'
$ws.Cells.Item(35,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(35,5).Value = 3
$ws.Cells.Item(35,6).Value = 3.5
# Row 36: id 4942597
$ws.Cells.Item(36,1).NumberFormat = "@"
$ws.Cells.Item(36,1).Value = '4942597'
$ws.Cells.Item(36,2).Value = 'Skross Adapter Europa-Schweiz'
$ws.Cells.Item(36,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/skross-adapter-europa-schweiz/p/4942597'
$ws.Cells.Item(36,4).ClearContents()
$ws.Cells.Item(36,7).Value = 'Skross'
$ws.Cells.Item(36,8).NumberFormat = "@"
$ws.Cells.Item(36,8).Value = '14.95'
$ws.Cells.Item(36,9).ClearContents()
$ws.Cells.Item(36,10).ClearContents()
$ws.Cells.Item(36,11).ClearContents()
$ws.Cells.Item(36,12).ClearContents()
$ws.Cells.Item(36,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(36,14).Value = 'Skross Adapter Europa-Schweiz 14.95 Schweizer Franken'
$ws.Cells.Item(36,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(36,5).ClearContents()
$ws.Cells.Item(36,6).Value = 0
# Row 37: id 3494065
$ws.Cells.Item(37,1).NumberFormat = "@"
$ws.Cells.Item(37,1).Value = '3494065'
$ws.Cells.Item(37,2).Value = 'Varta Longlife Max Power AA 4er Bli'
$ws.Cells.Item(37,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-max-power-aa-4er-bli/p/3494065'
$ws.Cells.Item(37,4).Value = '4ST'
$ws.Cells.Item(37,7).Value = 'Varta'
$ws.Cells.Item(37,8).NumberFormat = "@"
$ws.Cells.Item(37,8).Value = '9.95'
$ws.Cells.Item(37,9).Value = '2.49/1ST'
$ws.Cells.Item(37,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(37,11).NumberFormat = "@"
$ws.Cells.Item(37,11).Value = '2.49'
$ws.Cells.Item(37,12).Value = '1ST'
$ws.Cells.Item(37,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Cells.Item(37,14).Value = 'Varta Longlife Max Power AA 4er Bli 9.95 Schweizer Franken'
$ws.Cells.Item(37,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(37,5).Value = 1
$ws.Cells.Item(37,6).Value = 5
# Row 38: id 3494138
$ws.Cells.Item(38,1).NumberFormat = "@"
$ws.Cells.Item(38,1).Value = '3494138'
$ws.Cells.Item(38,2).Value = 'Varta Longlife Power Batterien AAA/LR03 8 Stück'
$ws.Cells.Item(38,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-power-batterien-aaalr03-8-stueck/p/3494138'
$ws.Cells.Item(38,4).Value = '8ST'
$ws.Cells.Item(38,7).Value = 'Varta'
$ws.Cells.Item(38,8).NumberFormat = "@"
$ws.Cells.Item(38,8).Value = '14.95'
$ws.Cells.Item(38,9).Value = '1.87/1ST'
$ws.Cells.Item(38,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(38,11).NumberFormat = "@"
$ws.Cells.Item(38,11).Value = '1.87'
$ws.Cells.Item(38,12).Value = '1ST'
$ws.Cells.Item(38,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Cells.Item(38,14).Value = 'Varta Longlife Power Batterien AAA/LR03 8 Stück 14.95 Schweizer Franken'
$ws.Cells.Item(38,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(38,5).Value = 1
$ws.Cells.Item(38,6).Value = 4
# Row 39: id 6753555
$ws.Cells.Item(39,1).NumberFormat = "@"
$ws.Cells.Item(39,1).Value = '6753555'
$ws.Cells.Item(39,2).Value = 'Duracell Batterien Optimum AAA/LR03 4 Stück'
$ws.Cells.Item(39,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/duracell-batterien-optimum-aaalr03-4-stueck/p/6753555'
$ws.Cells.Item(39,4).Value = '4ST'
$ws.Cells.Item(39,7).Value = 'Duracell'
$ws.Cells.Item(39,8).NumberFormat = "@"
$ws.Cells.Item(39,8).Value = '11.95'
$ws.Cells.Item(39,9).Value = '2.99/1ST'
$ws.Cells.Item(39,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(39,11).NumberFormat = "@"
$ws.Cells.Item(39,11).Value = '2.99'
$ws.Cells.Item(39,12).Value = '1ST'
$ws.Cells.Item(39,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Cells.Item(39,14).Value = 'Duracell Batterien Optimum AAA/LR03 4 Stück 11.95 Schweizer Franken'
$ws.Cells.Item(39,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(39,5).ClearContents()
$ws.Cells.Item(39,6).Value = 0
# Row 40: id 3494909
$ws.Cells.Item(40,1).NumberFormat = "@"
$ws.Cells.Item(40,1).Value = '3494909'
$ws.Cells.Item(40,2).Value = 'Varta Electronics V23GA 1er Bli'
$ws.Cells.Item(40,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-v23ga-1er-bli/p/3494909'
$ws.Cells.Item(40,4).Value = '1ST'
$ws.Cells.Item(40,7).Value = 'Varta'
$ws.Cells.Item(40,8).NumberFormat = "@"
$ws.Cells.Item(40,8).Value = '3.95'
$ws.Cells.Item(40,9).Value = '3.95/1ST'
$ws.Cells.Item(40,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(40,11).NumberFormat = "@"
$ws.Cells.Item(40,11).Value = '3.95'
$ws.Cells.Item(40,12).Value = '1ST'
$ws.Cells.Item(40,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(40,14).Value = 'Varta Electronics V23GA 1er Bli 3.95 Schweizer Franken'
$ws.Cells.Item(40,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(40,5).Value = 1
$ws.Cells.Item(40,6).Value = 5
# Row 41: id 4881282
$ws.Cells.Item(41,1).NumberFormat = "@"
$ws.Cells.Item(41,1).Value = '4881282'
$ws.Cells.Item(41,2).Value = 'Miele FJM HyClean 3D (4 Stück)'
$ws.Cells.Item(41,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/miele-fjm-hyclean-3d-4-stueck/p/4881282'
$ws.Cells.Item(41,4).ClearContents()
$ws.Cells.Item(41,7).Value = 'Miele'
$ws.Cells.Item(41,8).NumberFormat = "@"
$ws.Cells.Item(41,8).Value = '8.95'
$ws.Cells.Item(41,9).ClearContents()
$ws.Cells.Item(41,10).ClearContents()
$ws.Cells.Item(41,11).ClearContents()
$ws.Cells.Item(41,12).ClearContents()
$ws.Cells.Item(41,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(41,14).Value = 'Miele FJM HyClean 3D (4 Stück) 50% Aktion 8.95 Schweizer Franken statt 17.95 Schweizer Franken'
$ws.Cells.Item(41,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(41,5).ClearContents()
$ws.Cells.Item(41,6).Value = 0
# Row 42: id 4096751
$ws.Cells.Item(42,1).NumberFormat = "@"
$ws.Cells.Item(42,1).Value = '4096751'
$ws.Cells.Item(42,2).Value = 'Varta Longlife Power Batterien AA/LR6 6 Stück'
$ws.Cells.Item(42,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-power-batterien-aalr6-6-stueck/p/4096751'
$ws.Cells.Item(42,4).Value = '6ST'
$ws.Cells.Item(42,7).Value = 'Varta'
$ws.Cells.Item(42,8).NumberFormat = "@"
$ws.Cells.Item(42,8).Value = '12.95'
$ws.Cells.Item(42,9).Value = '2.16/1ST'
$ws.Cells.Item(42,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(42,11).NumberFormat = "@"
$ws.Cells.Item(42,11).Value = '2.16'
$ws.Cells.Item(42,12).Value = '1ST'
$ws.Cells.Item(42,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Cells.Item(42,14).Value = 'Varta Longlife Power Batterien AA/LR6 6 Stück 12.95 Schweizer Franken'
$ws.Cells.Item(42,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(42,5).ClearContents()
$ws.Cells.Item(42,6).Value = 0
# Row 43: id 4358323
$ws.Cells.Item(43,1).NumberFormat = "@"
$ws.Cells.Item(43,1).Value = '4358323'
$ws.Cells.Item(43,2).Value = 'Rayovac Hörgerätebatterien 312 6 Stück'
$ws.Cells.Item(43,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/rayovac-hoergeraetebatterien-312-6-stueck/p/4358323'
$ws.Cells.Item(43,4).Value = '6ST'
$ws.Cells.Item(43,7).Value = 'Rayovac'
$ws.Cells.Item(43,8).NumberFormat = "@"
$ws.Cells.Item(43,8).Value = '9.95'
$ws.Cells.Item(43,9).Value = '1.66/1ST'
$ws.Cells.Item(43,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(43,11).NumberFormat = "@"
$ws.Cells.Item(43,11).Value = '1.66'
$ws.Cells.Item(43,12).Value = '1ST'
$ws.Cells.Item(43,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(43,14).Value = 'Rayovac Hörgerätebatterien 312 6 Stück 9.95 Schweizer Franken'
$ws.Cells.Item(43,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(43,5).Value = 1
$ws.Cells.Item(43,6).Value = 4
# Row 44: id 5751576
$ws.Cells.Item(44,1).NumberFormat = "@"
$ws.Cells.Item(44,1).Value = '5751576'
$ws.Cells.Item(44,2).Value = 'satrap Toasty 1 Toaster'
$ws.Cells.Item(44,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-toasty-1-toaster/p/5751576'
$ws.Cells.Item(44,4).ClearContents()
$ws.Cells.Item(44,7).Value = 'satrap'
$ws.Cells.Item(44,8).NumberFormat = "@"
$ws.Cells.Item(44,8).Value = '22.45'
$ws.Cells.Item(44,9).ClearContents()
$ws.Cells.Item(44,10).ClearContents()
$ws.Cells.Item(44,11).ClearContents()
$ws.Cells.Item(44,12).ClearContents()
$ws.Cells.Item(44,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''kuechengeraete'']'
$ws.Cells.Item(44,14).Value = 'satrap Toasty 1 Toaster 25% Aktion 22.45 Schweizer Franken statt 29.95 Schweizer Franken'
$ws.Cells.Item(44,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(44,5).Value = 1
$ws.Cells.Item(44,6).Value = 5
# Row 45: id 3494063
$ws.Cells.Item(45,1).NumberFormat = "@"
$ws.Cells.Item(45,1).Value = '3494063'
$ws.Cells.Item(45,2).Value = 'Varta Longlife Power C 2er Bli'
$ws.Cells.Item(45,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-power-c-2er-bli/p/3494063'
$ws.Cells.Item(45,4).Value = '2ST'
$ws.Cells.Item(45,7).Value = 'Varta'
$ws.Cells.Item(45,8).NumberFormat = "@"
$ws.Cells.Item(45,8).Value = '7.95'
$ws.Cells.Item(45,9).Value = '3.98/1ST'
$ws.Cells.Item(45,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(45,11).NumberFormat = "@"
$ws.Cells.Item(45,11).Value = '3.98'
$ws.Cells.Item(45,12).Value = '1ST'
$ws.Cells.Item(45,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(45,14).Value = 'Varta Longlife Power C 2er Bli 7.95 Schweizer Franken'
$ws.Cells.Item(45,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(45,5).ClearContents()
$ws.Cells.Item(45,6).Value = 0
# Row 46: id 3494007
$ws.Cells.Item(46,1).NumberFormat = "@"
$ws.Cells.Item(46,1).Value = '3494007'
$ws.Cells.Item(46,2).Value = 'Varta Longlife Power D 2er Bli'
$ws.Cells.Item(46,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-power-d-2er-bli/p/3494007'
$ws.Cells.Item(46,4).Value = '2ST'
$ws.Cells.Item(46,7).Value = 'Varta'
$ws.Cells.Item(46,8).NumberFormat = "@"
$ws.Cells.Item(46,8).Value = '8.95'
$ws.Cells.Item(46,9).Value = '4.48/1ST'
$ws.Cells.Item(46,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(46,11).NumberFormat = "@"
$ws.Cells.Item(46,11).Value = '4.48'
$ws.Cells.Item(46,12).Value = '1ST'
$ws.Cells.Item(46,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(46,14).Value = 'Varta Longlife Power D 2er Bli 8.95 Schweizer Franken'
$ws.Cells.Item(46,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(46,5).ClearContents()
$ws.Cells.Item(46,6).Value = 0
# Row 47: id 6266870
$ws.Cells.Item(47,1).NumberFormat = "@"
$ws.Cells.Item(47,1).Value = '6266870'
$ws.Cells.Item(47,2).Value = 'Reer Steckdosen-Sicherung Kappe Schweiz 8 Stück'
$ws.Cells.Item(47,3).Value = '/de/haushalt-tier/elektroartikel-batterien/reer-steckdosen-sicherung-kappe-schweiz-8-stueck/p/6266870'
$ws.Cells.Item(47,4).ClearContents()
$ws.Cells.Item(47,7).Value = 'reer'
$ws.Cells.Item(47,8).NumberFormat = "@"
$ws.Cells.Item(47,8).Value = '3.95'
$ws.Cells.Item(47,9).ClearContents()
$ws.Cells.Item(47,10).ClearContents()
$ws.Cells.Item(47,11).ClearContents()
$ws.Cells.Item(47,12).ClearContents()
$ws.Cells.Item(47,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'']'
$ws.Cells.Item(47,14).Value = 'Reer Steckdosen-Sicherung Kappe Schweiz 8 Stück 3.95 Schweizer Franken'
$ws.Cells.Item(47,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(47,5).ClearContents()
$ws.Cells.Item(47,6).Value = 0
# Row 48: id 6119284
$ws.Cells.Item(48,1).NumberFormat = "@"
$ws.Cells.Item(48,1).Value = '6119284'
$ws.Cells.Item(48,2).Value = 'satrap Aqua SA10 Wasserkocher'
$ws.Cells.Item(48,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-aqua-sa10-wasserkocher/p/6119284'
$ws.Cells.Item(48,4).ClearContents()
$ws.Cells.Item(48,7).Value = 'satrap'
$ws.Cells.Item(48,8).NumberFormat = "@"
$ws.Cells.Item(48,8).Value = '37.45'
$ws.Cells.Item(48,9).ClearContents()
$ws.Cells.Item(48,10).ClearContents()
$ws.Cells.Item(48,11).ClearContents()
$ws.Cells.Item(48,12).ClearContents()
$ws.Cells.Item(48,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''kuechengeraete'']'
$ws.Cells.Item(48,14).Value = 'satrap Aqua SA10 Wasserkocher 25% Aktion 37.45 Schweizer Franken statt 49.95 Schweizer Franken'
$ws.Cells.Item(48,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(48,5).Value = 3
$ws.Cells.Item(48,6).Value = 2.5
# Row 49: id 6670192
$ws.Cells.Item(49,1).NumberFormat = "@"
$ws.Cells.Item(49,1).Value = '6670192'
$ws.Cells.Item(49,2).Value = 'Philips Wasserkocher HD9318/01'
$ws.Cells.Item(49,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/philips-wasserkocher-hd931801/p/6670192'
$ws.Cells.Item(49,4).ClearContents()
$ws.Cells.Item(49,7).Value = 'Philips'
$ws.Cells.Item(49,8).NumberFormat = "@"
$ws.Cells.Item(49,8).Value = '24.95'
$ws.Cells.Item(49,9).ClearContents()
$ws.Cells.Item(49,10).ClearContents()
$ws.Cells.Item(49,11).ClearContents()
$ws.Cells.Item(49,12).ClearContents()
$ws.Cells.Item(49,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''kuechengeraete'']'
$ws.Cells.Item(49,14).Value = 'Philips Wasserkocher HD9318/01 50% Aktion 24.95 Schweizer Franken statt 49.90 Schweizer Franken'
$ws.Cells.Item(49,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(49,5).ClearContents()
$ws.Cells.Item(49,6).Value = 0
# Row 50: id 5763068
$ws.Cells.Item(50,1).NumberFormat = "@"
$ws.Cells.Item(50,1).Value = '5763068'
$ws.Cells.Item(50,2).Value = 'satrap Aspira Sine A700 Beutelloser Zyklonstaubsauger'
$ws.Cells.Item(50,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-aspira-sine-a700-beutelloser-zyklonstaubsauger/p/5763068'
$ws.Cells.Item(50,4).ClearContents()
$ws.Cells.Item(50,7).Value = 'satrap'
$ws.Cells.Item(50,8).NumberFormat = "@"
$ws.Cells.Item(50,8).Value = '89.25'
$ws.Cells.Item(50,9).ClearContents()
$ws.Cells.Item(50,10).ClearContents()
$ws.Cells.Item(50,11).ClearContents()
$ws.Cells.Item(50,12).ClearContents()
$ws.Cells.Item(50,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(50,14).Value = 'satrap Aspira Sine A700 Beutelloser Zyklonstaubsauger 25% Aktion 89.25 Schweizer Franken statt 119.00 Schweizer Franken'
$ws.Cells.Item(50,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(50,5).Value = 4
$ws.Cells.Item(50,6).Value = 3
# Row 51: id 6298048
$ws.Cells.Item(51,1).NumberFormat = "@"
$ws.Cells.Item(51,1).Value = '6298048'
$ws.Cells.Item(51,2).Value = 'Electrolux E201S Classic Staubbeutel'
$ws.Cells.Item(51,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/electrolux-e201s-classic-staubbeutel/p/6298048'
$ws.Cells.Item(51,4).ClearContents()
$ws.Cells.Item(51,7).Value = 'Electrolux'
$ws.Cells.Item(51,8).NumberFormat = "@"
$ws.Cells.Item(51,8).Value = '9.95'
$ws.Cells.Item(51,9).ClearContents()
$ws.Cells.Item(51,10).ClearContents()
$ws.Cells.Item(51,11).ClearContents()
$ws.Cells.Item(51,12).ClearContents()
$ws.Cells.Item(51,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(51,14).Value = 'Electrolux E201S Classic Staubbeutel 9.95 Schweizer Franken'
$ws.Cells.Item(51,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(51,5).ClearContents()
$ws.Cells.Item(51,6).Value = 0
# Row 52: id 4942594
$ws.Cells.Item(52,1).NumberFormat = "@"
$ws.Cells.Item(52,1).Value = '4942594'
$ws.Cells.Item(52,2).Value = 'Skross Adapter World-Schweiz'
$ws.Cells.Item(52,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/skross-adapter-world-schweiz/p/4942594'
$ws.Cells.Item(52,4).ClearContents()
$ws.Cells.Item(52,7).Value = 'Skross'
$ws.Cells.Item(52,8).NumberFormat = "@"
$ws.Cells.Item(52,8).Value = '19.95'
$ws.Cells.Item(52,9).ClearContents()
$ws.Cells.Item(52,10).ClearContents()
$ws.Cells.Item(52,11).ClearContents()
$ws.Cells.Item(52,12).ClearContents()
$ws.Cells.Item(52,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(52,14).Value = 'Skross Adapter World-Schweiz 19.95 Schweizer Franken'
$ws.Cells.Item(52,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(52,5).ClearContents()
$ws.Cells.Item(52,6).Value = 0
# Row 53: id 5872158
$ws.Cells.Item(53,1).NumberFormat = "@"
$ws.Cells.Item(53,1).Value = '5872158'
$ws.Cells.Item(53,2).Value = 'satrap Tischventilator Venti 1'
$ws.Cells.Item(53,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-tischventilator-venti-1/p/5872158'
$ws.Cells.Item(53,4).ClearContents()
$ws.Cells.Item(53,7).Value = 'satrap'
$ws.Cells.Item(53,8).NumberFormat = "@"
$ws.Cells.Item(53,8).Value = '22.45'
$ws.Cells.Item(53,9).ClearContents()
$ws.Cells.Item(53,10).ClearContents()
$ws.Cells.Item(53,11).ClearContents()
$ws.Cells.Item(53,12).ClearContents()
$ws.Cells.Item(53,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(53,14).Value = 'satrap Tischventilator Venti 1 25% Aktion 22.45 Schweizer Franken statt 29.95 Schweizer Franken'
$ws.Cells.Item(53,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(53,5).ClearContents()
$ws.Cells.Item(53,6).Value = 0
# Row 54: id 5872164
$ws.Cells.Item(54,1).NumberFormat = "@"
$ws.Cells.Item(54,1).Value = '5872164'
$ws.Cells.Item(54,2).Value = 'Satrap Tischventilator Venti 2'
$ws.Cells.Item(54,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-tischventilator-venti-2/p/5872164'
$ws.Cells.Item(54,4).ClearContents()
$ws.Cells.Item(54,7).Value = 'satrap'
$ws.Cells.Item(54,8).NumberFormat = "@"
$ws.Cells.Item(54,8).Value = '26.20'
$ws.Cells.Item(54,9).ClearContents()
$ws.Cells.Item(54,10).ClearContents()
$ws.Cells.Item(54,11).ClearContents()
$ws.Cells.Item(54,12).ClearContents()
$ws.Cells.Item(54,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(54,14).Value = 'Satrap Tischventilator Venti 2 25% Aktion 26.20 Schweizer Franken statt 34.95 Schweizer Franken'
$ws.Cells.Item(54,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(54,5).ClearContents()
$ws.Cells.Item(54,6).Value = 0
# Row 55: id 5750424
$ws.Cells.Item(55,1).NumberFormat = "@"
$ws.Cells.Item(55,1).Value = '5750424'
$ws.Cells.Item(55,2).Value = 'satrap Vento HT2 Reisehaartrockner'
$ws.Cells.Item(55,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/beautygeraete/satrap-vento-ht2-reisehaartrockner/p/5750424'
$ws.Cells.Item(55,4).ClearContents()
$ws.Cells.Item(55,7).Value = 'satrap'
$ws.Cells.Item(55,8).NumberFormat = "@"
$ws.Cells.Item(55,8).Value = '22.45'
$ws.Cells.Item(55,9).ClearContents()
$ws.Cells.Item(55,10).ClearContents()
$ws.Cells.Item(55,11).ClearContents()
$ws.Cells.Item(55,12).ClearContents()
$ws.Cells.Item(55,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''beautygeraete'']'
$ws.Cells.Item(55,14).Value = 'satrap Vento HT2 Reisehaartrockner 25% Aktion 22.45 Schweizer Franken statt 29.95 Schweizer Franken'
$ws.Cells.Item(55,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(55,5).Value = 1
$ws.Cells.Item(55,6).Value = 5
# Row 56: id 4905486
$ws.Cells.Item(56,1).NumberFormat = "@"
$ws.Cells.Item(56,1).Value = '4905486'
$ws.Cells.Item(56,2).Value = 'Alkaline Batterie 3LR12/4.5V'
$ws.Cells.Item(56,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/alkaline-batterie-3lr1245v/p/4905486'
$ws.Cells.Item(56,4).Value = '1ST'
$ws.Cells.Item(56,7).Value = 'Coop'
$ws.Cells.Item(56,8).NumberFormat = "@"
$ws.Cells.Item(56,8).Value = '5.95'
$ws.Cells.Item(56,9).Value = '5.95/1ST'
$ws.Cells.Item(56,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(56,11).NumberFormat = "@"
$ws.Cells.Item(56,11).Value = '5.95'
$ws.Cells.Item(56,12).Value = '1ST'
$ws.Cells.Item(56,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(56,14).Value = 'Alkaline Batterie 3LR12/4.5V 5.95 Schweizer Franken'
$ws.Cells.Item(56,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(56,5).Value = 1
$ws.Cells.Item(56,6).Value = 1
# Row 57: id 3494231
$ws.Cells.Item(57,1).NumberFormat = "@"
$ws.Cells.Item(57,1).Value = '3494231'
$ws.Cells.Item(57,2).Value = 'Varta Electronics CR2025 1er Bli'
$ws.Cells.Item(57,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-cr2025-1er-bli/p/3494231'
$ws.Cells.Item(57,4).Value = '1ST'
$ws.Cells.Item(57,7).Value = 'Varta'
$ws.Cells.Item(57,8).NumberFormat = "@"
$ws.Cells.Item(57,8).Value = '4.95'
$ws.Cells.Item(57,9).Value = '4.95/1ST'
$ws.Cells.Item(57,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(57,11).NumberFormat = "@"
$ws.Cells.Item(57,11).Value = '4.95'
$ws.Cells.Item(57,12).Value = '1ST'
$ws.Cells.Item(57,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(57,14).Value = 'Varta Electronics CR2025 1er Bli 4.95 Schweizer Franken'
$ws.Cells.Item(57,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(57,5).ClearContents()
$ws.Cells.Item(57,6).Value = 0
# Row 58: id 6973029
$ws.Cells.Item(58,1).NumberFormat = "@"
$ws.Cells.Item(58,1).Value = '6973029'
$ws.Cells.Item(58,2).Value = 'Severin Tischgrill PG 8565'
$ws.Cells.Item(58,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/severin-tischgrill-pg-8565/p/6973029'
$ws.Cells.Item(58,4).ClearContents()
$ws.Cells.Item(58,7).Value = 'Severin'
$ws.Cells.Item(58,8).NumberFormat = "@"
$ws.Cells.Item(58,8).Value = '49.95'
$ws.Cells.Item(58,9).ClearContents()
$ws.Cells.Item(58,10).ClearContents()
$ws.Cells.Item(58,11).ClearContents()
$ws.Cells.Item(58,12).ClearContents()
$ws.Cells.Item(58,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''kuechengeraete'']'
$ws.Cells.Item(58,14).Value = 'Severin Tischgrill PG 8565 50% Aktion 49.95 Schweizer Franken statt 99.95 Schweizer Franken'
$ws.Cells.Item(58,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(58,5).Value = 1
$ws.Cells.Item(58,6).Value = 3
# Row 59: id 4096628
$ws.Cells.Item(59,1).NumberFormat = "@"
$ws.Cells.Item(59,1).Value = '4096628'
$ws.Cells.Item(59,2).Value = 'Varta Longlife Power Batterien AAA/LR03 6 Stück'
$ws.Cells.Item(59,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-power-batterien-aaalr03-6-stueck/p/4096628'
$ws.Cells.Item(59,4).Value = '6ST'
$ws.Cells.Item(59,7).Value = 'Varta'
$ws.Cells.Item(59,8).NumberFormat = "@"
$ws.Cells.Item(59,8).Value = '12.95'
$ws.Cells.Item(59,9).Value = '2.16/1ST'
$ws.Cells.Item(59,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(59,11).NumberFormat = "@"
$ws.Cells.Item(59,11).Value = '2.16'
$ws.Cells.Item(59,12).Value = '1ST'
$ws.Cells.Item(59,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(59,14).Value = 'Varta Longlife Power Batterien AAA/LR03 6 Stück 12.95 Schweizer Franken'
$ws.Cells.Item(59,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(59,5).ClearContents()
$ws.Cells.Item(59,6).Value = 0
# Row 60: id 4589934
$ws.Cells.Item(60,1).NumberFormat = "@"
$ws.Cells.Item(60,1).Value = '4589934'
$ws.Cells.Item(60,2).Value = 'Varta Longlife AA 4er Bli'
$ws.Cells.Item(60,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-longlife-aa-4er-bli/p/4589934'
$ws.Cells.Item(60,4).Value = '4ST'
$ws.Cells.Item(60,7).Value = 'Varta'
$ws.Cells.Item(60,8).NumberFormat = "@"
$ws.Cells.Item(60,8).Value = '7.95'
$ws.Cells.Item(60,9).Value = '1.99/1ST'
$ws.Cells.Item(60,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(60,11).NumberFormat = "@"
$ws.Cells.Item(60,11).Value = '1.99'
$ws.Cells.Item(60,12).Value = '1ST'
$ws.Cells.Item(60,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Cells.Item(60,14).Value = 'Varta Longlife AA 4er Bli 7.95 Schweizer Franken'
$ws.Cells.Item(60,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(60,5).ClearContents()
$ws.Cells.Item(60,6).Value = 0
# Row 61: id 6805078
$ws.Cells.Item(61,1).NumberFormat = "@"
$ws.Cells.Item(61,1).Value = '6805078'
$ws.Cells.Item(61,2).Value = 'Varta Premium Light F20'
$ws.Cells.Item(61,3).Value = '/de/haushalt-tier/elektroartikel-batterien/beleuchtung/taschenlampen-lichter/varta-premium-light-f20/p/6805078'
$ws.Cells.Item(61,4).ClearContents()
$ws.Cells.Item(61,7).Value = 'Varta'
$ws.Cells.Item(61,8).NumberFormat = "@"
$ws.Cells.Item(61,8).Value = '39.95'
$ws.Cells.Item(61,9).ClearContents()
$ws.Cells.Item(61,10).ClearContents()
$ws.Cells.Item(61,11).ClearContents()
$ws.Cells.Item(61,12).ClearContents()
$ws.Cells.Item(61,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''beleuchtung'', ''taschenlampen-lichter'']'
$ws.Cells.Item(61,14).Value = 'Varta Premium Light F20 39.95 Schweizer Franken'
$ws.Cells.Item(61,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(61,5).ClearContents()
$ws.Cells.Item(61,6).Value = 0
# Row 62: id 6378875
$ws.Cells.Item(62,1).NumberFormat = "@"
$ws.Cells.Item(62,1).Value = '6378875'
$ws.Cells.Item(62,2).Value = 'Wilkinson Intuition 4in1 Finish Styler'
$ws.Cells.Item(62,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/beautygeraete/wilkinson-intuition-4in1-finish-styler/p/6378875'
$ws.Cells.Item(62,4).Value = '1ST'
$ws.Cells.Item(62,7).Value = 'Wilkinson'
$ws.Cells.Item(62,8).NumberFormat = "@"
$ws.Cells.Item(62,8).Value = '29.95'
$ws.Cells.Item(62,9).Value = '29.95/1ST'
$ws.Cells.Item(62,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(62,11).NumberFormat = "@"
$ws.Cells.Item(62,11).Value = '29.95'
$ws.Cells.Item(62,12).Value = '1ST'
$ws.Cells.Item(62,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''beautygeraete'']'
$ws.Cells.Item(62,14).Value = 'Wilkinson Intuition 4in1 Finish Styler 29.95 Schweizer Franken'
$ws.Cells.Item(62,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(62,5).ClearContents()
$ws.Cells.Item(62,6).Value = 0
# Row 63: id 6579165
$ws.Cells.Item(63,1).NumberFormat = "@"
$ws.Cells.Item(63,1).Value = '6579165'
$ws.Cells.Item(63,2).Value = 'satrap Trimm Set XA 5-in-1 Haar/Trimm-Set'
$ws.Cells.Item(63,3).Value = '/de/kosmetik-gesundheit/herrenpflege-rasur/elektrische-herrenrasierer/satrap-trimm-set-xa-5-in-1-haartrimm-set/p/6579165'
$ws.Cells.Item(63,4).ClearContents()
$ws.Cells.Item(63,7).Value = 'satrap'
$ws.Cells.Item(63,8).NumberFormat = "@"
$ws.Cells.Item(63,8).Value = '44.95'
$ws.Cells.Item(63,9).ClearContents()
$ws.Cells.Item(63,10).ClearContents()
$ws.Cells.Item(63,11).ClearContents()
$ws.Cells.Item(63,12).ClearContents()
$ws.Cells.Item(63,13).Value = '[''kosmetik-gesundheit'', ''herrenpflege-rasur'', ''elektrische-herrenrasierer'']'
$ws.Cells.Item(63,14).Value = 'satrap Trimm Set XA 5-in-1 Haar/Trimm-Set 25% Aktion 44.95 Schweizer Franken statt 59.95 Schweizer Franken'
$ws.Cells.Item(63,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(63,5).ClearContents()
$ws.Cells.Item(63,6).Value = 0
# Row 64: id 3494230
$ws.Cells.Item(64,1).NumberFormat = "@"
$ws.Cells.Item(64,1).Value = '3494230'
$ws.Cells.Item(64,2).Value = 'Varta Electronics V13GS / V357 1er Bli'
$ws.Cells.Item(64,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-electronics-v13gs--v357-1er-bli/p/3494230'
$ws.Cells.Item(64,4).Value = '1ST'
$ws.Cells.Item(64,7).Value = 'Varta'
$ws.Cells.Item(64,8).NumberFormat = "@"
$ws.Cells.Item(64,8).Value = '4.95'
$ws.Cells.Item(64,9).Value = '4.95/1ST'
$ws.Cells.Item(64,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(64,11).NumberFormat = "@"
$ws.Cells.Item(64,11).Value = '4.95'
$ws.Cells.Item(64,12).Value = '1ST'
$ws.Cells.Item(64,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(64,14).Value = 'Varta Electronics V13GS / V357 1er Bli 4.95 Schweizer Franken'
$ws.Cells.Item(64,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(64,5).ClearContents()
$ws.Cells.Item(64,6).Value = 0
# Row 65: id 4119046
$ws.Cells.Item(65,1).NumberFormat = "@"
$ws.Cells.Item(65,1).Value = '4119046'
$ws.Cells.Item(65,2).Value = 'Varta Ultra Lithium AA 4er Bli'
$ws.Cells.Item(65,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-ultra-lithium-aa-4er-bli/p/4119046'
$ws.Cells.Item(65,4).Value = '4ST'
$ws.Cells.Item(65,7).Value = 'Varta'
$ws.Cells.Item(65,8).NumberFormat = "@"
$ws.Cells.Item(65,8).Value = '14.95'
$ws.Cells.Item(65,9).Value = '3.74/1ST'
$ws.Cells.Item(65,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(65,11).NumberFormat = "@"
$ws.Cells.Item(65,11).Value = '3.74'
$ws.Cells.Item(65,12).Value = '1ST'
$ws.Cells.Item(65,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aa'']'
$ws.Cells.Item(65,14).Value = 'Varta Ultra Lithium AA 4er Bli 14.95 Schweizer Franken'
$ws.Cells.Item(65,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(65,5).Value = 1
$ws.Cells.Item(65,6).Value = 5
# Row 66: id 6700107
$ws.Cells.Item(66,1).NumberFormat = "@"
$ws.Cells.Item(66,1).Value = '6700107'
$ws.Cells.Item(66,2).Value = 'MIELE Compact C2 Classic (mit Beutel, 890 W)'
$ws.Cells.Item(66,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/miele-compact-c2-classic-mit-beutel-890-w/p/6700107'
$ws.Cells.Item(66,4).ClearContents()
$ws.Cells.Item(66,7).Value = 'Miele'
$ws.Cells.Item(66,8).NumberFormat = "@"
$ws.Cells.Item(66,8).Value = '194.50'
$ws.Cells.Item(66,9).ClearContents()
$ws.Cells.Item(66,10).ClearContents()
$ws.Cells.Item(66,11).ClearContents()
$ws.Cells.Item(66,12).ClearContents()
$ws.Cells.Item(66,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(66,14).Value = 'MIELE Compact C2 Classic (mit Beutel, 890 W) 50% Aktion 194.50 Schweizer Franken statt 389.00 Schweizer Franken'
$ws.Cells.Item(66,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(66,5).ClearContents()
$ws.Cells.Item(66,6).Value = 0
# Row 67: id 6007534
$ws.Cells.Item(67,1).NumberFormat = "@"
$ws.Cells.Item(67,1).Value = '6007534'
$ws.Cells.Item(67,2).Value = 'Trend USB-Stick 8 GB'
$ws.Cells.Item(67,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-8-gb/p/6007534'
$ws.Cells.Item(67,4).ClearContents()
$ws.Cells.Item(67,7).Value = 'Trend'
$ws.Cells.Item(67,8).NumberFormat = "@"
$ws.Cells.Item(67,8).Value = '9.95'
$ws.Cells.Item(67,9).ClearContents()
$ws.Cells.Item(67,10).ClearContents()
$ws.Cells.Item(67,11).ClearContents()
$ws.Cells.Item(67,12).ClearContents()
$ws.Cells.Item(67,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(67,14).Value = 'Trend USB-Stick 8 GB 9.95 Schweizer Franken'
$ws.Cells.Item(67,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(67,5).Value = 1
$ws.Cells.Item(67,6).Value = 5
# Row 68: id 4425412
$ws.Cells.Item(68,1).NumberFormat = "@"
$ws.Cells.Item(68,1).Value = '4425412'
$ws.Cells.Item(68,2).Value = 'Militärtaschenlampe LED'
$ws.Cells.Item(68,3).Value = '/de/haushalt-tier/elektroartikel-batterien/beleuchtung/taschenlampen-lichter/militaertaschenlampe-led/p/4425412'
$ws.Cells.Item(68,4).ClearContents()
$ws.Cells.Item(68,7).Value = 'Coop'
$ws.Cells.Item(68,8).NumberFormat = "@"
$ws.Cells.Item(68,8).Value = '8.95'
$ws.Cells.Item(68,9).ClearContents()
$ws.Cells.Item(68,10).ClearContents()
$ws.Cells.Item(68,11).ClearContents()
$ws.Cells.Item(68,12).ClearContents()
$ws.Cells.Item(68,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''beleuchtung'', ''taschenlampen-lichter'']'
$ws.Cells.Item(68,14).Value = 'Militärtaschenlampe LED 8.95 Schweizer Franken'
$ws.Cells.Item(68,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(68,5).ClearContents()
$ws.Cells.Item(68,6).Value = 0
# Row 69: id 3591269
$ws.Cells.Item(69,1).NumberFormat = "@"
$ws.Cells.Item(69,1).Value = '3591269'
$ws.Cells.Item(69,2).Value = 'Varta Longlife Max Power C 2er Bli'
$ws.Cells.Item(69,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-max-power-c-2er-bli/p/3591269'
$ws.Cells.Item(69,4).Value = '2ST'
$ws.Cells.Item(69,7).Value = 'Varta'
$ws.Cells.Item(69,8).NumberFormat = "@"
$ws.Cells.Item(69,8).Value = '8.95'
$ws.Cells.Item(69,9).Value = '4.48/1ST'
$ws.Cells.Item(69,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(69,11).NumberFormat = "@"
$ws.Cells.Item(69,11).Value = '4.48'
$ws.Cells.Item(69,12).Value = '1ST'
$ws.Cells.Item(69,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(69,14).Value = 'Varta Longlife Max Power C 2er Bli 8.95 Schweizer Franken'
$ws.Cells.Item(69,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(69,5).Value = 1
$ws.Cells.Item(69,6).Value = 5
# Row 70: id 5751209
$ws.Cells.Item(70,1).NumberFormat = "@"
$ws.Cells.Item(70,1).Value = '5751209'
$ws.Cells.Item(70,2).Value = 'satrap Ferru BE1 2200 Dampfbügeleisen'
$ws.Cells.Item(70,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-ferru-be1-2200-dampfbuegeleisen/p/5751209'
$ws.Cells.Item(70,4).ClearContents()
$ws.Cells.Item(70,7).Value = 'satrap'
$ws.Cells.Item(70,8).NumberFormat = "@"
$ws.Cells.Item(70,8).Value = '37.45'
$ws.Cells.Item(70,9).ClearContents()
$ws.Cells.Item(70,10).ClearContents()
$ws.Cells.Item(70,11).ClearContents()
$ws.Cells.Item(70,12).ClearContents()
$ws.Cells.Item(70,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(70,14).Value = 'satrap Ferru BE1 2200 Dampfbügeleisen 25% Aktion 37.45 Schweizer Franken statt 49.95 Schweizer Franken'
$ws.Cells.Item(70,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(70,5).Value = 3
$ws.Cells.Item(70,6).Value = 2.5
# Row 71: id 4358322
$ws.Cells.Item(71,1).NumberFormat = "@"
$ws.Cells.Item(71,1).Value = '4358322'
$ws.Cells.Item(71,2).Value = 'Rayovac Hörgerätebatterien 13 6 Stück'
$ws.Cells.Item(71,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/rayovac-hoergeraetebatterien-13-6-stueck/p/4358322'
$ws.Cells.Item(71,4).Value = '6ST'
$ws.Cells.Item(71,7).Value = 'Rayovac'
$ws.Cells.Item(71,8).NumberFormat = "@"
$ws.Cells.Item(71,8).Value = '9.95'
$ws.Cells.Item(71,9).Value = '1.66/1ST'
$ws.Cells.Item(71,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(71,11).NumberFormat = "@"
$ws.Cells.Item(71,11).Value = '1.66'
$ws.Cells.Item(71,12).Value = '1ST'
$ws.Cells.Item(71,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(71,14).Value = 'Rayovac Hörgerätebatterien 13 6 Stück 9.95 Schweizer Franken'
$ws.Cells.Item(71,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(71,5).Value = 1
$ws.Cells.Item(71,6).Value = 5
# Row 72: id 6425996
$ws.Cells.Item(72,1).NumberFormat = "@"
$ws.Cells.Item(72,1).Value = '6425996'
$ws.Cells.Item(72,2).Value = 'satrap espresso XA Kolbenkaffeemaschine'
$ws.Cells.Item(72,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-espresso-xa-kolbenkaffeemaschine/p/6425996'
$ws.Cells.Item(72,4).ClearContents()
$ws.Cells.Item(72,7).Value = 'satrap'
$ws.Cells.Item(72,8).NumberFormat = "@"
$ws.Cells.Item(72,8).Value = '186.75'
$ws.Cells.Item(72,9).ClearContents()
$ws.Cells.Item(72,10).ClearContents()
$ws.Cells.Item(72,11).ClearContents()
$ws.Cells.Item(72,12).ClearContents()
$ws.Cells.Item(72,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''kuechengeraete'']'
$ws.Cells.Item(72,14).Value = 'satrap espresso XA Kolbenkaffeemaschine 25% Aktion 186.75 Schweizer Franken statt 249.00 Schweizer Franken'
$ws.Cells.Item(72,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(72,5).Value = 7
$ws.Cells.Item(72,6).Value = 4
# Row 73: id 6725106
$ws.Cells.Item(73,1).NumberFormat = "@"
$ws.Cells.Item(73,1).Value = '6725106'
$ws.Cells.Item(73,2).Value = 'Braun Dampfbügelstation IS1012BL'
$ws.Cells.Item(73,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/braun-dampfbuegelstation-is1012bl/p/6725106'
$ws.Cells.Item(73,4).ClearContents()
$ws.Cells.Item(73,7).Value = 'Braun'
$ws.Cells.Item(73,8).NumberFormat = "@"
$ws.Cells.Item(73,8).Value = '149.00'
$ws.Cells.Item(73,9).ClearContents()
$ws.Cells.Item(73,10).ClearContents()
$ws.Cells.Item(73,11).ClearContents()
$ws.Cells.Item(73,12).ClearContents()
$ws.Cells.Item(73,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(73,14).Value = 'Braun Dampfbügelstation IS1012BL 149.00 Schweizer Franken'
$ws.Cells.Item(73,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(73,5).ClearContents()
$ws.Cells.Item(73,6).Value = 0
# Row 74: id 6995204
$ws.Cells.Item(74,1).NumberFormat = "@"
$ws.Cells.Item(74,1).Value = '6995204'
$ws.Cells.Item(74,2).Value = 'Electrolux AirFryer Range Explore 6'
$ws.Cells.Item(74,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/electrolux-airfryer-range-explore-6/p/6995204'
$ws.Cells.Item(74,4).ClearContents()
$ws.Cells.Item(74,7).Value = 'Electrolux'
$ws.Cells.Item(74,8).NumberFormat = "@"
$ws.Cells.Item(74,8).Value = '199.00'
$ws.Cells.Item(74,9).ClearContents()
$ws.Cells.Item(74,10).ClearContents()
$ws.Cells.Item(74,11).ClearContents()
$ws.Cells.Item(74,12).ClearContents()
$ws.Cells.Item(74,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''kuechengeraete'']'
$ws.Cells.Item(74,14).Value = 'Electrolux AirFryer Range Explore 6 199.00 Schweizer Franken'
$ws.Cells.Item(74,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(74,5).ClearContents()
$ws.Cells.Item(74,6).Value = 0
# Row 75: id 6735643
$ws.Cells.Item(75,1).NumberFormat = "@"
$ws.Cells.Item(75,1).Value = '6735643'
$ws.Cells.Item(75,2).Value = 'LED 31V Anschlussset Transf.+Verl.kabel'
$ws.Cells.Item(75,3).Value = '/de/haushalt-tier/haushalt-kueche/uebrige-haushaltsartikel/led-31v-anschlussset-transfverlkabel/p/6735643'
$ws.Cells.Item(75,4).ClearContents()
$ws.Cells.Item(75,7).Value = 'Coop'
$ws.Cells.Item(75,8).NumberFormat = "@"
$ws.Cells.Item(75,8).Value = '9.95'
$ws.Cells.Item(75,9).ClearContents()
$ws.Cells.Item(75,10).ClearContents()
$ws.Cells.Item(75,11).ClearContents()
$ws.Cells.Item(75,12).ClearContents()
$ws.Cells.Item(75,13).Value = '[''haushalt-tier'', ''haushalt-kueche'', ''uebrige-haushaltsartikel'']'
$ws.Cells.Item(75,14).Value = 'LED 31V Anschlussset Transf.+Verl.kabel 50% Aktion 9.95 Schweizer Franken statt 19.95 Schweizer Franken'
$ws.Cells.Item(75,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(75,5).Value = 1
$ws.Cells.Item(75,6).Value = 5
# Row 76: id 5882124
$ws.Cells.Item(76,1).NumberFormat = "@"
$ws.Cells.Item(76,1).Value = '5882124'
$ws.Cells.Item(76,2).Value = 'Philips Avent Audio Monitors DECT-Babyphone'
$ws.Cells.Item(76,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/philips-avent-audio-monitors-dect-babyphone/p/5882124'
$ws.Cells.Item(76,4).ClearContents()
$ws.Cells.Item(76,7).Value = 'Avent'
$ws.Cells.Item(76,8).NumberFormat = "@"
$ws.Cells.Item(76,8).Value = '99.90'
$ws.Cells.Item(76,9).ClearContents()
$ws.Cells.Item(76,10).ClearContents()
$ws.Cells.Item(76,11).ClearContents()
$ws.Cells.Item(76,12).ClearContents()
$ws.Cells.Item(76,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'']'
$ws.Cells.Item(76,14).Value = 'Philips Avent Audio Monitors DECT-Babyphone 99.90 Schweizer Franken'
$ws.Cells.Item(76,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(76,5).ClearContents()
$ws.Cells.Item(76,6).Value = 0
# Row 77: id 6689620
$ws.Cells.Item(77,1).NumberFormat = "@"
$ws.Cells.Item(77,1).Value = '6689620'
$ws.Cells.Item(77,2).Value = 'Philips Bügeleisen GC1750/21'
$ws.Cells.Item(77,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/philips-buegeleisen-gc175021/p/6689620'
$ws.Cells.Item(77,4).ClearContents()
$ws.Cells.Item(77,7).Value = 'Philips'
$ws.Cells.Item(77,8).NumberFormat = "@"
$ws.Cells.Item(77,8).Value = '59.90'
$ws.Cells.Item(77,9).ClearContents()
$ws.Cells.Item(77,10).ClearContents()
$ws.Cells.Item(77,11).ClearContents()
$ws.Cells.Item(77,12).ClearContents()
$ws.Cells.Item(77,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(77,14).Value = 'Philips Bügeleisen GC1750/21 59.90 Schweizer Franken'
$ws.Cells.Item(77,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(77,5).ClearContents()
$ws.Cells.Item(77,6).Value = 0
# Row 78: id 5764836
$ws.Cells.Item(78,1).NumberFormat = "@"
$ws.Cells.Item(78,1).Value = '5764836'
$ws.Cells.Item(78,2).Value = 'satrap Aqua W4 Wasserkocher'
$ws.Cells.Item(78,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-aqua-w4-wasserkocher/p/5764836'
$ws.Cells.Item(78,4).ClearContents()
$ws.Cells.Item(78,7).Value = 'satrap'
$ws.Cells.Item(78,8).NumberFormat = "@"
$ws.Cells.Item(78,8).Value = '29.95'
$ws.Cells.Item(78,9).ClearContents()
$ws.Cells.Item(78,10).ClearContents()
$ws.Cells.Item(78,11).ClearContents()
$ws.Cells.Item(78,12).ClearContents()
$ws.Cells.Item(78,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''kuechengeraete'']'
$ws.Cells.Item(78,14).Value = 'satrap Aqua W4 Wasserkocher 25% Aktion 29.95 Schweizer Franken statt 39.95 Schweizer Franken'
$ws.Cells.Item(78,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(78,5).ClearContents()
$ws.Cells.Item(78,6).Value = 0
# Row 79: id 6288875
$ws.Cells.Item(79,1).NumberFormat = "@"
$ws.Cells.Item(79,1).Value = '6288875'
$ws.Cells.Item(79,2).Value = 'satrap Bügeleisen mit Power-Dampfstation Station XA 6.5'
$ws.Cells.Item(79,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/satrap-buegeleisen-mit-power-dampfstation-station-xa-65/p/6288875'
$ws.Cells.Item(79,4).ClearContents()
$ws.Cells.Item(79,7).Value = 'satrap'
$ws.Cells.Item(79,8).NumberFormat = "@"
$ws.Cells.Item(79,8).Value = '111.75'
$ws.Cells.Item(79,9).ClearContents()
$ws.Cells.Item(79,10).ClearContents()
$ws.Cells.Item(79,11).ClearContents()
$ws.Cells.Item(79,12).ClearContents()
$ws.Cells.Item(79,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(79,14).Value = 'satrap Bügeleisen mit Power-Dampfstation Station XA 6.5 25% Aktion 111.75 Schweizer Franken statt 149.00 Schweizer Franken'
$ws.Cells.Item(79,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(79,5).Value = 1
$ws.Cells.Item(79,6).Value = 2
# Row 80: id 6125818
$ws.Cells.Item(80,1).NumberFormat = "@"
$ws.Cells.Item(80,1).Value = '6125818'
$ws.Cells.Item(80,2).Value = 'satrap Mano XA Handmixer'
$ws.Cells.Item(80,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/satrap-mano-xa-handmixer/p/6125818'
$ws.Cells.Item(80,4).ClearContents()
$ws.Cells.Item(80,7).Value = 'satrap'
$ws.Cells.Item(80,8).NumberFormat = "@"
$ws.Cells.Item(80,8).Value = '37.45'
$ws.Cells.Item(80,9).ClearContents()
$ws.Cells.Item(80,10).ClearContents()
$ws.Cells.Item(80,11).ClearContents()
$ws.Cells.Item(80,12).ClearContents()
$ws.Cells.Item(80,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''kuechengeraete'']'
$ws.Cells.Item(80,14).Value = 'satrap Mano XA Handmixer 25% Aktion 37.45 Schweizer Franken statt 49.95 Schweizer Franken'
$ws.Cells.Item(80,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(80,5).Value = 1
$ws.Cells.Item(80,6).Value = 1
# Row 81: id 6753975
$ws.Cells.Item(81,1).NumberFormat = "@"
$ws.Cells.Item(81,1).Value = '6753975'
$ws.Cells.Item(81,2).Value = 'Severin Standgrill mit Grillplatte PG 8563'
$ws.Cells.Item(81,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/kuechengeraete/severin-standgrill-mit-grillplatte-pg-8563/p/6753975'
$ws.Cells.Item(81,4).ClearContents()
$ws.Cells.Item(81,7).Value = 'Severin'
$ws.Cells.Item(81,8).NumberFormat = "@"
$ws.Cells.Item(81,8).Value = '74.50'
$ws.Cells.Item(81,9).ClearContents()
$ws.Cells.Item(81,10).ClearContents()
$ws.Cells.Item(81,11).ClearContents()
$ws.Cells.Item(81,12).ClearContents()
$ws.Cells.Item(81,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''kuechengeraete'']'
$ws.Cells.Item(81,14).Value = 'Severin Standgrill mit Grillplatte PG 8563 50% Aktion 74.50 Schweizer Franken statt 149.00 Schweizer Franken'
$ws.Cells.Item(81,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(81,5).ClearContents()
$ws.Cells.Item(81,6).Value = 0
# Row 82: id 6822783
$ws.Cells.Item(82,1).NumberFormat = "@"
$ws.Cells.Item(82,1).Value = '6822783'
$ws.Cells.Item(82,2).Value = 'Trend Car Charger QC'
$ws.Cells.Item(82,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-car-charger-qc/p/6822783'
$ws.Cells.Item(82,4).ClearContents()
$ws.Cells.Item(82,7).Value = 'Trend'
$ws.Cells.Item(82,8).NumberFormat = "@"
$ws.Cells.Item(82,8).Value = '19.95'
$ws.Cells.Item(82,9).ClearContents()
$ws.Cells.Item(82,10).ClearContents()
$ws.Cells.Item(82,11).ClearContents()
$ws.Cells.Item(82,12).ClearContents()
$ws.Cells.Item(82,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(82,14).Value = 'Trend Car Charger QC 19.95 Schweizer Franken'
$ws.Cells.Item(82,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(82,5).ClearContents()
$ws.Cells.Item(82,6).Value = 0
# Row 83: id 6872922
$ws.Cells.Item(83,1).NumberFormat = "@"
$ws.Cells.Item(83,1).Value = '6872922'
$ws.Cells.Item(83,2).Value = 'Trend Radio DAB+ DB-23C'
$ws.Cells.Item(83,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/trend-radio-dab-db-23c/p/6872922'
$ws.Cells.Item(83,4).ClearContents()
$ws.Cells.Item(83,7).Value = 'Trend'
$ws.Cells.Item(83,8).NumberFormat = "@"
$ws.Cells.Item(83,8).Value = '49.95'
$ws.Cells.Item(83,9).ClearContents()
$ws.Cells.Item(83,10).ClearContents()
$ws.Cells.Item(83,11).ClearContents()
$ws.Cells.Item(83,12).ClearContents()
$ws.Cells.Item(83,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'']'
$ws.Cells.Item(83,14).Value = 'Trend Radio DAB+ DB-23C 49.95 Schweizer Franken'
$ws.Cells.Item(83,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(83,5).ClearContents()
$ws.Cells.Item(83,6).Value = 0
# Row 84: id 6007537
$ws.Cells.Item(84,1).NumberFormat = "@"
$ws.Cells.Item(84,1).Value = '6007537'
$ws.Cells.Item(84,2).Value = 'Trend USB-Stick 16 GB'
$ws.Cells.Item(84,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-16-gb/p/6007537'
$ws.Cells.Item(84,4).ClearContents()
$ws.Cells.Item(84,7).Value = 'Trend'
$ws.Cells.Item(84,8).NumberFormat = "@"
$ws.Cells.Item(84,8).Value = '16.95'
$ws.Cells.Item(84,9).ClearContents()
$ws.Cells.Item(84,10).ClearContents()
$ws.Cells.Item(84,11).ClearContents()
$ws.Cells.Item(84,12).ClearContents()
$ws.Cells.Item(84,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(84,14).Value = 'Trend USB-Stick 16 GB 16.95 Schweizer Franken'
$ws.Cells.Item(84,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(84,5).Value = 1
$ws.Cells.Item(84,6).Value = 2
# Row 85: id 6459240
$ws.Cells.Item(85,1).NumberFormat = "@"
$ws.Cells.Item(85,1).Value = '6459240'
$ws.Cells.Item(85,2).Value = 'Trend USB-Stick 256 GB'
$ws.Cells.Item(85,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-256-gb/p/6459240'
$ws.Cells.Item(85,4).ClearContents()
$ws.Cells.Item(85,7).Value = 'Trend'
$ws.Cells.Item(85,8).NumberFormat = "@"
$ws.Cells.Item(85,8).Value = '59.95'
$ws.Cells.Item(85,9).ClearContents()
$ws.Cells.Item(85,10).ClearContents()
$ws.Cells.Item(85,11).ClearContents()
$ws.Cells.Item(85,12).ClearContents()
$ws.Cells.Item(85,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(85,14).Value = 'Trend USB-Stick 256 GB 59.95 Schweizer Franken'
$ws.Cells.Item(85,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(85,5).ClearContents()
$ws.Cells.Item(85,6).Value = 0
# Row 86: id 6007538
$ws.Cells.Item(86,1).NumberFormat = "@"
$ws.Cells.Item(86,1).Value = '6007538'
$ws.Cells.Item(86,2).Value = 'Trend USB-Stick 32 GB'
$ws.Cells.Item(86,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-32-gb/p/6007538'
$ws.Cells.Item(86,4).ClearContents()
$ws.Cells.Item(86,7).Value = 'Trend'
$ws.Cells.Item(86,8).NumberFormat = "@"
$ws.Cells.Item(86,8).Value = '29.95'
$ws.Cells.Item(86,9).ClearContents()
$ws.Cells.Item(86,10).ClearContents()
$ws.Cells.Item(86,11).ClearContents()
$ws.Cells.Item(86,12).ClearContents()
$ws.Cells.Item(86,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(86,14).Value = 'Trend USB-Stick 32 GB 29.95 Schweizer Franken'
$ws.Cells.Item(86,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(86,5).ClearContents()
$ws.Cells.Item(86,6).Value = 0
# Row 87: id 5867973
$ws.Cells.Item(87,1).NumberFormat = "@"
$ws.Cells.Item(87,1).Value = '5867973'
$ws.Cells.Item(87,2).Value = 'Trend USB-Stick 64 GB'
$ws.Cells.Item(87,3).Value = '/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/haushaltsgeraete-kabel/trend-usb-stick-64-gb/p/5867973'
$ws.Cells.Item(87,4).ClearContents()
$ws.Cells.Item(87,7).Value = 'Trend'
$ws.Cells.Item(87,8).NumberFormat = "@"
$ws.Cells.Item(87,8).Value = '39.95'
$ws.Cells.Item(87,9).ClearContents()
$ws.Cells.Item(87,10).ClearContents()
$ws.Cells.Item(87,11).ClearContents()
$ws.Cells.Item(87,12).ClearContents()
$ws.Cells.Item(87,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''elektrogeraete'', ''haushaltsgeraete-kabel'']'
$ws.Cells.Item(87,14).Value = 'Trend USB-Stick 64 GB 39.95 Schweizer Franken'
$ws.Cells.Item(87,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(87,5).Value = 4
$ws.Cells.Item(87,6).Value = 2
# Row 88: id 5894674
$ws.Cells.Item(88,1).NumberFormat = "@"
$ws.Cells.Item(88,1).Value = '5894674'
$ws.Cells.Item(88,2).Value = 'Trisa Beauty Sonic Nail Care System'
$ws.Cells.Item(88,3).Value = '/de/kosmetik-gesundheit/make-up/naegel/set-accessoires/trisa-beauty-sonic-nail-care-system/p/5894674'
$ws.Cells.Item(88,4).Value = '1ST'
$ws.Cells.Item(88,7).Value = 'Trisa'
$ws.Cells.Item(88,8).NumberFormat = "@"
$ws.Cells.Item(88,8).Value = '29.95'
$ws.Cells.Item(88,9).Value = '29.95/1ST'
$ws.Cells.Item(88,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(88,11).NumberFormat = "@"
$ws.Cells.Item(88,11).Value = '29.95'
$ws.Cells.Item(88,12).Value = '1ST'
$ws.Cells.Item(88,13).Value = '[''kosmetik-gesundheit'', ''make-up'', ''naegel'', ''set-accessoires'']'
$ws.Cells.Item(88,14).Value = 'Trisa Beauty Sonic Nail Care System 29.95 Schweizer Franken'
$ws.Cells.Item(88,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(88,5).ClearContents()
$ws.Cells.Item(88,6).Value = 0
# Row 89: id 3862219
$ws.Cells.Item(89,1).NumberFormat = "@"
$ws.Cells.Item(89,1).Value = '3862219'
$ws.Cells.Item(89,2).Value = 'Varta Active LED Camping Lanterne'
$ws.Cells.Item(89,3).Value = '/de/haushalt-tier/elektroartikel-batterien/beleuchtung/taschenlampen-lichter/varta-active-led-camping-lanterne/p/3862219'
$ws.Cells.Item(89,4).ClearContents()
$ws.Cells.Item(89,7).Value = 'Varta'
$ws.Cells.Item(89,8).NumberFormat = "@"
$ws.Cells.Item(89,8).Value = '38.50'
$ws.Cells.Item(89,9).ClearContents()
$ws.Cells.Item(89,10).ClearContents()
$ws.Cells.Item(89,11).ClearContents()
$ws.Cells.Item(89,12).ClearContents()
$ws.Cells.Item(89,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''beleuchtung'', ''taschenlampen-lichter'']'
$ws.Cells.Item(89,14).Value = 'Varta Active LED Camping Lanterne 38.50 Schweizer Franken'
$ws.Cells.Item(89,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(89,5).ClearContents()
$ws.Cells.Item(89,6).Value = 0
# Row 90: id 4589933
$ws.Cells.Item(90,1).NumberFormat = "@"
$ws.Cells.Item(90,1).Value = '4589933'
$ws.Cells.Item(90,2).Value = 'Varta Longlife AAA 4er Bli'
$ws.Cells.Item(90,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-aaa-4er-bli/p/4589933'
$ws.Cells.Item(90,4).Value = '4ST'
$ws.Cells.Item(90,7).Value = 'Varta'
$ws.Cells.Item(90,8).NumberFormat = "@"
$ws.Cells.Item(90,8).Value = '7.95'
$ws.Cells.Item(90,9).Value = '1.99/1ST'
$ws.Cells.Item(90,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(90,11).NumberFormat = "@"
$ws.Cells.Item(90,11).Value = '1.99'
$ws.Cells.Item(90,12).Value = '1ST'
$ws.Cells.Item(90,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Cells.Item(90,14).Value = 'Varta Longlife AAA 4er Bli 7.95 Schweizer Franken'
$ws.Cells.Item(90,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(90,5).ClearContents()
$ws.Cells.Item(90,6).Value = 0
# Row 91: id 4589935
$ws.Cells.Item(91,1).NumberFormat = "@"
$ws.Cells.Item(91,1).Value = '4589935'
$ws.Cells.Item(91,2).Value = 'Varta Longlife C 2er Bli'
$ws.Cells.Item(91,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-c-2er-bli/p/4589935'
$ws.Cells.Item(91,4).Value = '2ST'
$ws.Cells.Item(91,7).Value = 'Varta'
$ws.Cells.Item(91,8).NumberFormat = "@"
$ws.Cells.Item(91,8).Value = '6.95'
$ws.Cells.Item(91,9).Value = '3.48/1ST'
$ws.Cells.Item(91,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(91,11).NumberFormat = "@"
$ws.Cells.Item(91,11).Value = '3.48'
$ws.Cells.Item(91,12).Value = '1ST'
$ws.Cells.Item(91,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(91,14).Value = 'Varta Longlife C 2er Bli 6.95 Schweizer Franken'
$ws.Cells.Item(91,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(91,5).ClearContents()
$ws.Cells.Item(91,6).Value = 0
# Row 92: id 3494067
$ws.Cells.Item(92,1).NumberFormat = "@"
$ws.Cells.Item(92,1).Value = '3494067'
$ws.Cells.Item(92,2).Value = 'Varta Longlife Max Power AAA 4er Bli'
$ws.Cells.Item(92,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/aaa/varta-longlife-max-power-aaa-4er-bli/p/3494067'
$ws.Cells.Item(92,4).Value = '4ST'
$ws.Cells.Item(92,7).Value = 'Varta'
$ws.Cells.Item(92,8).NumberFormat = "@"
$ws.Cells.Item(92,8).Value = '9.95'
$ws.Cells.Item(92,9).Value = '2.49/1ST'
$ws.Cells.Item(92,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(92,11).NumberFormat = "@"
$ws.Cells.Item(92,11).Value = '2.49'
$ws.Cells.Item(92,12).Value = '1ST'
$ws.Cells.Item(92,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''aaa'']'
$ws.Cells.Item(92,14).Value = 'Varta Longlife Max Power AAA 4er Bli 9.95 Schweizer Franken'
$ws.Cells.Item(92,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(92,5).ClearContents()
$ws.Cells.Item(92,6).Value = 0
# Row 93: id 3591272
$ws.Cells.Item(93,1).NumberFormat = "@"
$ws.Cells.Item(93,1).Value = '3591272'
$ws.Cells.Item(93,2).Value = 'Varta Longlife Max Power D 2er Bli'
$ws.Cells.Item(93,3).Value = '/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-max-power-d-2er-bli/p/3591272'
$ws.Cells.Item(93,4).Value = '2ST'
$ws.Cells.Item(93,7).Value = 'Varta'
$ws.Cells.Item(93,8).NumberFormat = "@"
$ws.Cells.Item(93,8).Value = '9.95'
$ws.Cells.Item(93,9).Value = '4.98/1ST'
$ws.Cells.Item(93,10).Value = 'Preis pro 1 Stück'
$ws.Cells.Item(93,11).NumberFormat = "@"
$ws.Cells.Item(93,11).Value = '4.98'
$ws.Cells.Item(93,12).Value = '1ST'
$ws.Cells.Item(93,13).Value = '[''haushalt-tier'', ''elektroartikel-batterien'', ''batterien'', ''andere-batterien'']'
$ws.Cells.Item(93,14).Value = 'Varta Longlife Max Power D 2er Bli 9.95 Schweizer Franken'
$ws.Cells.Item(93,15).Value = '2022-09-10 20:58:20'
$ws.Cells.Item(93,5).ClearContents()
$ws.Cells.Item(93,6).Value = 0

# Drop the now-redundant trailing row (old row 94 duplicate) so the
# used range shrinks back to A1:O93, matching the refreshed crawl output.
$ws.Rows.Item(94).Delete()
